{"js": "const replacements = [\n  { oldText: \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\u6975\u5ea6\u4f9d\u8cf4\u767c\u660e\u3002\u9019\u662f\u4eba\u985e\u5275\u9020\u6027\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\u5176\u6700\u7d42\u76ee\u7684\u5728\u65bc\u5b8c\u5168\u638c\u63a7\u5fc3\u9748\u5c0d\u7269\u8cea\u4e16\u754c\u7684\u652f\u914d\uff0c\u5c07\u81ea\u7136\u7684\u529b\u91cf\u99d5\u99ad\u65bc\u4eba\u985e\u7684\u9700\u6c42\u3002\u9019\u662f\u767c\u660e\u5bb6\u7684\u8271\u5de8\u4efb\u52d9\uff0c\u4ed6\u5e38\u88ab\u8aa4\u89e3\u4e14\u4e0d\u53d7\u91cd\u8996\u3002\u4f46\u4ed6\u5f9e\u5176\u80fd\u529b\u7684\u6109\u5feb\u904b\u7528\u4e2d\uff0c\u4ee5\u53ca\u5f9e\u77e5\u9053\u81ea\u5df1\u5c6c\u65bc\u4e00\u500b\u6c92\u6709\u4e86\u5c31\u4eba\u985e\u65e9\u5c31\u6703\u5728\u8207\u7121\u60c5\u5143\u7d20\u7684\u82e6\u82e6\u9b25\u722d\u4e2d\u6ec5\u7d55\u7684\u6975\u5176\u7279\u6b0a\u968e\u7d1a\u4e2d\uff0c\u627e\u5230\u4e86\u5145\u8db3\u7684\u88dc\u511f\u3002\", newText: \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\uff0c\u81f3\u95dc\u91cd\u8981\u5730\u4f9d\u8cf4\u65bc\u767c\u660e\u3002\u5b83\u662f\u4eba\u985e\u5275\u9020\u6027\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\u5176\u6700\u7d42\u76ee\u7684\u662f\u5fc3\u9748\u5c0d\u7269\u8cea\u4e16\u754c\u7684\u5b8c\u5168\u638c\u63e1\uff0c\u5c07\u81ea\u7136\u529b\u91cf\u99d5\u99ad\u65bc\u4eba\u985e\u9700\u6c42\u3002\u9019\u662f\u767c\u660e\u5bb6\u8271\u5de8\u7684\u4efb\u52d9\uff0c\u4ed6\u7d93\u5e38\u88ab\u8aa4\u89e3\u548c\u5f97\u4e0d\u5230\u5831\u916c\u3002\u4f46\u4ed6\u537b\u5728\u6109\u6085\u5730\u904b\u7528\u5176\u80fd\u529b\u7684\u904e\u7a0b\u4e2d\uff0c\u4ee5\u53ca\u5728\u610f\u8b58\u5230\u81ea\u5df1\u662f\u90a3\u4e00\u7fa4\u7279\u6b8a\u800c\u6709\u7279\u6b0a\u7684\u968e\u7d1a\u4e2d\u7684\u4e00\u54e1\u6642\uff0c\u7372\u5f97\u4e86\u5145\u5206\u7684\u88dc\u511f\uff0c\u6c92\u6709\u4ed6\u5011\uff0c\u4eba\u985e\u7a2e\u65cf\u65e9\u5c31\u6703\u5728\u8207\u7121\u60c5\u5143\u7d20\u7684\u82e6\u82e6\u9b25\u722d\u4e2d\u6ec5\u4ea1\u3002\" },\n  { oldText: \"\u5c31\u6211\u500b\u4eba\u800c\u8a00\uff0c\u6211\u5df2\u7d93\u5145\u5206\u9ad4\u9a57\u904e\u9019\u7a2e\u7121\u8207\u502b\u6bd4\u7684\u6a02\u8da3\uff0c\u751a\u81f3\u591a\u5e74\u4f86\u6211\u7684\u751f\u6d3b\u5e7e\u4e4e\u7b49\u540c\u65bc\u6301\u7e8c\u7684\u72c2\u559c\u3002\u6211\u88ab\u8a8d\u70ba\u662f\u6700\u52e4\u596e\u7684\u4eba\u4e4b\u4e00\uff0c\u4e5f\u8a31\u6211\u78ba\u5be6\u5982\u6b64\uff0c\u5982\u679c\u601d\u8003\u7b49\u540c\u65bc\u52de\u52d5\uff0c\u56e0\u70ba\u6211\u5e7e\u4e4e\u628a\u6240\u6709\u6e05\u9192\u7684\u6642\u9593\u90fd\u6295\u5165\u5176\u4e2d\u3002\u4f46\u5982\u679c\u5de5\u4f5c\u88ab\u89e3\u91cb\u70ba\u6839\u64da\u56b4\u683c\u898f\u5247\u5728\u7279\u5b9a\u6642\u9593\u5167\u5b8c\u6210\u7684\u660e\u78ba\u5de5\u4f5c\uff0c\u90a3\u9ebc\u6211\u53ef\u80fd\u662f\u6700\u61f6\u60f0\u7684\u4eba\u3002\u4efb\u4f55\u88ab\u8feb\u7684\u52aa\u529b\u90fd\u9700\u8981\u72a7\u7272\u751f\u547d\u80fd\u91cf\u3002\u6211\u5f9e\u672a\u4ed8\u51fa\u904e\u9019\u6a23\u7684\u4ee3\u50f9\u3002\u76f8\u53cd\u5730\uff0c\u6211\u5f9e\u601d\u60f3\u4e2d\u7372\u76ca\u532a\u6dfa\u3002\", newText: \"\u5c31\u6211\u500b\u4eba\u800c\u8a00\uff0c\u6211\u5df2\u7d93\u4eab\u53d7\u5230\u9019\u7a2e\u7cbe\u7dfb\u7684\u6a02\u8da3\uff0c\u800c\u4e14\u662f\u8d85\u904e\u98fd\u548c\u7684\u7a0b\u5ea6\uff0c\u4ee5\u81f3\u65bc\u591a\u5e74\u4f86\u6211\u7684\u751f\u6d3b\u5e7e\u4e4e\u7b49\u540c\u65bc\u6301\u7e8c\u7684\u72c2\u559c\u3002\u6211\u88ab\u8a8d\u70ba\u662f\u6700\u52aa\u529b\u5de5\u4f5c\u7684\u4eba\u4e4b\u4e00\uff0c\u4e5f\u8a31\u6211\u78ba\u5be6\u662f\uff0c\u5982\u679c\u628a\u601d\u8003\u7b49\u540c\u65bc\u52de\u52d5\uff0c\u56e0\u70ba\u6211\u5e7e\u4e4e\u5c07\u6240\u6709\u7684\u6e05\u9192\u6642\u9593\u90fd\u5949\u737b\u7d66\u4e86\u5b83\u3002\u4f46\u662f\uff0c\u5982\u679c\u5de5\u4f5c\u88ab\u7406\u89e3\u70ba\u5728\u898f\u5b9a\u7684\u6642\u9593\u5167\u6309\u7167\u56b4\u683c\u7684\u898f\u5247\u9032\u884c\u7684\u660e\u78ba\u8868\u73fe\uff0c\u90a3\u9ebc\u6211\u53ef\u80fd\u662f\u6700\u61f6\u60f0\u7684\u4eba\u3002\u4efb\u4f55\u88ab\u8feb\u7684\u52aa\u529b\u90fd\u8981\u6c42\u72a7\u7272\u751f\u547d\u80fd\u91cf\u3002\u6211\u5f9e\u672a\u4ed8\u51fa\u904e\u9019\u6a23\u7684\u4ee3\u50f9\u3002\u76f8\u53cd\uff0c\u6211\u5f9e\u81ea\u5df1\u7684\u601d\u60f3\u4e2d\u7372\u5f97\u4e86\u6ecb\u990a\u3002\" },\n  { oldText: \"\u5728\u8a66\u5716\u4ee5\u9023\u8cab\u4e14\u5fe0\u5be6\u7684\u65b9\u5f0f\u8a18\u9304\u6211\u5728\u9019\u7cfb\u5217\u6587\u7ae0\u4e2d\u7684\u6d3b\u52d5\u6642\uff0c\u9019\u4e9b\u6587\u7ae0\u5c07\u5728\u300a\u96fb\u6c23\u5be6\u9a57\u5bb6\u300b\u7de8\u8f2f\u7684\u5354\u52a9\u4e0b\u5448\u73fe\uff0c\u4e3b\u8981\u9762\u5411\u5e74\u8f15\u7537\u6027\u8b80\u8005\uff0c\u6211\u4e0d\u5f97\u4e0d\uff0c\u5118\u7ba1\u52c9\u5f37\uff0c\u56de\u9867\u6211\u9752\u6625\u6642\u671f\u7684\u5370\u8c61\uff0c\u4ee5\u53ca\u90a3\u4e9b\u6c7a\u5b9a\u6211\u8077\u696d\u751f\u6daf\u7684\u74b0\u5883\u548c\u4e8b\u4ef6\u3002\", newText: \"\u5728\u5617\u8a66\u4ee5\u9023\u8cab\u4e14\u5fe0\u5be6\u7684\u65b9\u5f0f\uff0c\u85c9\u7531\u300a\u96fb\u6c23\u5be6\u9a57\u5bb6\u300b\u7de8\u8f2f\u5011\u7684\u5354\u52a9\uff0c\u5728\u672c\u7cfb\u5217\u6587\u7ae0\u4e2d\u5411\u5e74\u8f15\u8b80\u8005\u5011\u5448\u73fe\u6211\u7684\u6d3b\u52d5\u8a18\u9304\u6642\uff0c\u6211\u5fc5\u9808\u4e0d\u60c5\u9858\u5730\u56de\u9867\u6211\u7684\u9752\u5c11\u5e74\u6642\u671f\uff0c\u4ee5\u53ca\u90a3\u4e9b\u6c7a\u5b9a\u6211\u8077\u696d\u751f\u6daf\u7684\u74b0\u5883\u548c\u4e8b\u4ef6\u3002\" },\n  { oldText: \"\u6211\u5011\u6700\u521d\u7684\u52aa\u529b\u7d14\u7cb9\u662f\u672c\u80fd\u7684\uff0c\u4f86\u81ea\u751f\u52d5\u800c\u4e0d\u53d7\u7d04\u675f\u7684\u60f3\u50cf\u529b\u7684\u9a45\u4f7f\u3002\u96a8\u8457\u5e74\u9f61\u7684\u589e\u9577\uff0c\u7406\u6027\u958b\u59cb\u4e3b\u5f35\u81ea\u5df1\u7684\u5730\u4f4d\uff0c\u6211\u5011\u8b8a\u5f97\u8d8a\u4f86\u8d8a\u6709\u689d\u7406\u548c\u8a2d\u8a08\u6027\u3002\u4f46\u90a3\u4e9b\u65e9\u671f\u7684\u885d\u52d5\uff0c\u96d6\u7136\u4e0d\u4e00\u5b9a\u7acb\u5373\u7522\u751f\u6210\u679c\uff0c\u537b\u662f\u6700\u91cd\u8981\u7684\uff0c\u53ef\u80fd\u6703\u5851\u9020\u6211\u5011\u7684\u4eba\u751f\u8ecc\u8de1\u3002\u4e8b\u5be6\u4e0a\uff0c\u6211\u73fe\u5728\u611f\u5230\uff0c\u5982\u679c\u6211\u7576\u6642\u7406\u89e3\u4e26\u57f9\u990a\u5b83\u5011\uff0c\u800c\u4e0d\u662f\u58d3\u6291\u5b83\u5011\uff0c\u6211\u5c07\u70ba\u4e16\u754c\u7559\u4e0b\u7684\u907a\u7522\u589e\u6dfb\u5de8\u5927\u7684\u50f9\u503c\u3002\u4f46\u662f\uff0c\u76f4\u5230\u6211\u6210\u5e74\u5f8c\u624d\u610f\u8b58\u5230\u81ea\u5df1\u662f\u4e00\u4f4d\u767c\u660e\u5bb6\u3002\", newText: \"\u6211\u5011\u6700\u521d\u7684\u52aa\u529b\u7d14\u7cb9\u662f\u672c\u80fd\u7684\uff0c\u662f\u9bae\u6d3b\u800c\u4e0d\u53d7\u7d04\u675f\u7684\u60f3\u50cf\u529b\u7684\u9a45\u4f7f\u3002\u96a8\u8457\u6211\u5011\u5e74\u9f61\u7684\u589e\u9577\uff0c\u7406\u6027\u958b\u59cb asserted \u81ea\u6211\uff0c\u6211\u5011\u8b8a\u5f97\u8d8a\u4f86\u8d8a\u7cfb\u7d71\u5316\u548c\u6709\u8a2d\u8a08\u6027\u3002\u4f46\u90a3\u4e9b\u65e9\u671f\u7684\u885d\u52d5\uff0c\u5118\u7ba1\u6c92\u6709\u7acb\u5373\u7522\u751f\u6210\u679c\uff0c\u537b\u662f\u81f3\u95dc\u91cd\u8981\u7684\uff0c\u5b83\u5011\u53ef\u80fd\u6703\u5851\u9020\u6211\u5011\u7684\u4eba\u751f\u8ecc\u8de1\u3002\u4e8b\u5be6\u4e0a\uff0c\u6211\u73fe\u5728\u610f\u8b58\u5230\uff0c\u5982\u679c\u6211\u80fd\u7406\u89e3\u4e26\u57f9\u990a\u5b83\u5011\uff0c\u800c\u4e0d\u662f\u58d3\u6291\u5b83\u5011\uff0c\u6211\u5c07\u70ba\u4e16\u754c\u7559\u4e0b\u7684\u907a\u7522\u589e\u6dfb\u91cd\u5927\u7684\u50f9\u503c\u3002\u4f46\u76f4\u5230\u6211\u6210\u5e74\u5f8c\uff0c\u6211\u624d\u610f\u8b58\u5230\u81ea\u5df1\u662f\u4e00\u500b\u767c\u660e\u5bb6\u3002\" },\n  { oldText: \"\u9019\u4e00\u5207\u6e90\u65bc\u591a\u7a2e\u539f\u56e0\u3002\u9996\u5148\uff0c\u6211\u6709\u4e00\u500b\u5929\u8ce6\u7570\u7a1f\u7684\u5144\u5f1f\u2014\u2014\u90a3\u662f\u751f\u7269\u5b78\u7814\u7a76\u7121\u6cd5\u89e3\u91cb\u7684\u7f55\u898b\u5fc3\u667a\u73fe\u8c61\u4e4b\u4e00\u3002\u4ed6\u904e\u65e9\u7684\u96e2\u4e16\u8b93\u7236\u6bcd\u60b2\u75db\u6b32\u7d55\u3002\u6211\u5011\u64c1\u6709\u4e00\u5339\u7531\u4e00\u4f4d\u646f\u53cb\u8d08\u9001\u7684\u99ac\u3002\u5b83\u662f\u4e00\u5339\u963f\u62c9\u4f2f\u8840\u7d71\u7684\u96c4\u58ef\u52d5\u7269\uff0c\u64c1\u6709\u8fd1\u4e4e\u4eba\u985e\u7684\u667a\u6167\uff0c\u53d7\u5230\u5168\u5bb6\u4eba\u5bf5\u611b\u548c\u5475\u8b77\uff0c\u751a\u81f3\u66fe\u5728\u4e00\u6b21\u975e\u51e1\u7684\u4e8b\u4ef6\u4e2d\u6551\u4e86\u7236\u89aa\u7684\u6027\u547d\u3002\u4e00\u500b\u51ac\u591c\uff0c\u7236\u89aa\u63a5\u5230\u7dca\u6025\u4efb\u52d9\uff0c\u5728\u7a7f\u8d8a\u72fc\u7fa4\u51fa\u6c92\u7684\u5c71\u5340\u6642\uff0c\u99ac\u5339\u53d7\u9a5a\u9003\u8dd1\uff0c\u5c07\u4ed6\u731b\u70c8\u5730\u6454\u5012\u5728\u5730\u3002\u99ac\u5339\u56de\u5230\u5bb6\u4e2d\uff0c\u6e3e\u8eab\u662f\u8840\uff0c\u7cbe\u75b2\u529b\u76e1\uff0c\u4f46\u8b66\u5831\u97ff\u8d77\u5f8c\uff0c\u5b83\u7acb\u5373\u518d\u6b21\u51fa\u767c\uff0c\u8fd4\u56de\u539f\u8655\uff0c\u641c\u6551\u968a\u9084\u6c92\u8d70\u591a\u9060\uff0c\u5c31\u9047\u5230\u4e86\u6062\u5fa9\u610f\u8b58\u4e26\u91cd\u65b0\u9a0e\u4e0a\u99ac\u7684\u7236\u89aa\uff0c\u4ed6\u4e26\u4e0d\u77e5\u9053\u81ea\u5df1\u5df2\u7d93\u5728\u96ea\u5730\u4e2d\u8eba\u4e86\u5e7e\u500b\u5c0f\u6642\u3002\u9019\u5339\u99ac\u5c0e\u81f4\u4e86\u5f1f\u5f1f\u7684\u53d7\u50b7\uff0c\u6700\u7d42\u5c0e\u81f4\u4ed6\u6b7b\u4ea1\u3002\u6211\u89aa\u773c\u898b\u8b49\u4e86\u9019\u5834\u60b2\u5287\uff0c\u5118\u7ba1\u5df2\u7d93\u904e\u4e86\u4e94\u5341\u516d\u5e74\uff0c\u4f46\u5c0d\u90a3\u5834\u666f\u7684\u8996\u89ba\u5370\u8c61\u4ecd\u7136\u6b77\u6b77\u5728\u76ee\u3002\u56de\u61b6\u8d77\u4ed6\u7684\u6210\u5c31\uff0c\u8b93\u6211\u6240\u6709\u7684\u52aa\u529b\u90fd\u986f\u5f97\u9eef\u7136\u5931\u8272\u3002\", newText: \"\u9019\u7de3\u65bc\u8a31\u591a\u539f\u56e0\u3002\u9996\u5148\uff0c\u6211\u6709\u4e00\u4f4d\u5929\u8ce6\u7570\u7a1f\u7684\u54e5\u54e5\u2014\u2014\u5c6c\u65bc\u751f\u7269\u5b78\u7814\u7a76\u7121\u6cd5\u89e3\u91cb\u7684\u7f55\u898b\u5fc3\u667a\u73fe\u8c61\u4e4b\u4e00\u3002\u4ed6\u65e9\u901d\u4f7f\u7236\u6bcd\u60b2\u75db\u6b32\u7d55\u3002\u6211\u5011\u64c1\u6709\u4e00\u5339\u7531\u4e00\u4f4d\u646f\u53cb\u8d08\u9001\u7684\u99ac\u3002\u9019\u662f\u4e00\u5339\u963f\u62c9\u4f2f\u54c1\u7a2e\u7684\u96c4\u5049\u52d5\u7269\uff0c\u5e7e\u4e4e\u5177\u6709\u4eba\u985e\u7684\u667a\u6167\uff0c\u53d7\u5230\u5168\u5bb6\u4eba\u7167\u6599\u548c\u5bf5\u611b\uff0c\u66fe\u5728\u7279\u6b8a\u60c5\u6cc1\u4e0b\u6551\u904e\u6211\u7236\u89aa\u7684\u547d\u3002\u6709\u4e00\u6b21\uff0c\u5728\u4e00\u500b\u51ac\u591c\uff0c\u6211\u7236\u89aa\u88ab\u7dca\u6025\u53ec\u559a\u53bb\u57f7\u884c\u4e00\u9805\u4efb\u52d9\uff0c\u7576\u4ed6\u7a7f\u8d8a\u72fc\u7fa4\u51fa\u6c92\u7684\u5c71\u5340\u6642\uff0c\u9019\u5339\u99ac\u53d7\u9a5a\u800c\u9003\uff0c\u5c07\u4ed6\u731b\u70c8\u5730\u62cb\u5230\u5730\u4e0a\u3002\u5b83\u56de\u5230\u5bb6\u6642\uff0c\u8840\u8de1\u6591\u6591\uff0c\u7cbe\u75b2\u529b\u7aed\uff0c\u4f46\u8b66\u5831\u97ff\u8d77\u5f8c\uff0c\u5b83\u7acb\u523b\u518d\u6b21\u885d\u51fa\u53bb\uff0c\u56de\u5230\u4e8b\u767c\u5730\u9ede\u3002\u5728\u641c\u5c0b\u968a\u51fa\u767c\u4e0d\u4e45\uff0c\u4ed6\u5011\u5c31\u9047\u5230\u4e86\u6211\u7236\u89aa\uff0c\u4ed6\u5df2\u6062\u5fa9\u610f\u8b58\uff0c\u91cd\u65b0\u9a0e\u4e0a\u99ac\uff0c\u4e26\u672a\u610f\u8b58\u5230\u81ea\u5df1\u5728\u96ea\u5730\u88e1\u8eba\u4e86\u5e7e\u500b\u5c0f\u6642\u3002\u9019\u5339\u99ac\u9020\u6210\u6211\u54e5\u54e5\u7684\u50b7\u52e2\uff0c\u5c0e\u81f4\u4ed6\u6b7b\u4ea1\u3002\u6211\u89aa\u773c\u76ee\u7779\u4e86\u9019\u5834\u60b2\u5287\uff0c\u96d6\u7136\u5df2\u7d93\u904e\u4e86\u4e94\u5341\u516d\u5e74\uff0c\u4f46\u6211\u5c0d\u5b83\u7684\u8996\u89ba\u5370\u8c61\u7d72\u6beb\u672a\u6e1b\u3002\u56de\u61b6\u8d77\u4ed6\u7684\u6210\u5c31\uff0c\u8b93\u6211\u6240\u6709\u7684\u52aa\u529b\u90fd\u986f\u5f97\u9eef\u7136\u5931\u8272\u3002\" },\n  { oldText: \"\u6211\u6240\u505a\u7684\u4e00\u5207\u503c\u5f97\u7a31\u8b9a\u7684\u4e8b\uff0c\u53ea\u6703\u8b93\u6211\u7684\u7236\u6bcd\u66f4\u52a0\u75db\u60dc\u4ed6\u5011\u7684\u640d\u5931\u3002\u56e0\u6b64\uff0c\u6211\u5728\u6210\u9577\u904e\u7a0b\u4e2d\u7f3a\u4e4f\u81ea\u4fe1\u3002\u4f46\u6211\u7d55\u975e\u611a\u8822\u7684\u7537\u5b69\uff0c\u81f3\u5c11\u5f9e\u6211\u81f3\u4eca\u8a18\u61b6\u7336\u65b0\u7684\u4e8b\u4ef6\u4e2d\u53ef\u4ee5\u770b\u51fa\u3002\u6709\u4e00\u5929\uff0c\u5e02\u8b70\u54e1\u5011\u8def\u904e\u4e00\u689d\u8857\uff0c\u6211\u6b63\u5728\u90a3\u88e1\u548c\u5176\u4ed6\u7684\u7537\u5b69\u73a9\u800d\u3002\u9019\u4e9b\u5c0a\u8cb4\u7684\u9577\u8005\u4e2d\uff0c\u4e00\u4f4d\u5bcc\u6709\u7684\u5e02\u6c11\u505c\u4e0b\u4f86\uff0c\u7d66\u4e86\u6211\u5011\u6bcf\u500b\u4eba\u4e00\u679a\u9280\u5e63\u3002\u7576\u4ed6\u8d70\u5230\u6211\u9762\u524d\u6642\uff0c\u7a81\u7136\u505c\u4f4f\uff0c\u547d\u4ee4\u9053\uff1a\u300c\u770b\u8457\u6211\u7684\u773c\u775b\u3002\u300d\u6211\u8207\u4ed6\u7684\u76ee\u5149\u76f8\u9047\uff0c\u4f38\u51fa\u624b\u53bb\u63a5\u90a3\u679a\u73cd\u8cb4\u7684\u786c\u5e63\uff0c\u7136\u800c\uff0c\u4ee4\u6211\u6cae\u55aa\u7684\u662f\uff0c\u4ed6\u8aaa\uff1a\u300c\u4e0d\uff0c\u4e0d\u7528\u4e86\uff0c\u4f60\u5f9e\u6211\u9019\u88e1\u5f97\u4e0d\u5230\u4efb\u4f55\u6771\u897f\uff0c\u4f60\u592a\u8070\u660e\u4e86\u3002\u300d\u4ed6\u5011\u7d93\u5e38\u8b1b\u4e00\u500b\u95dc\u65bc\u6211\u7684\u6709\u8da3\u6545\u4e8b\u3002\u6211\u6709\u5169\u500b\u81c9\u4e0a\u6709\u76ba\u7d0b\u7684\u8001\u963f\u59e8\uff0c\u5176\u4e2d\u4e00\u500b\u6709\u5169\u9846\u7259\u9f52\u50cf\u8c61\u7259\u4e00\u6a23\u7a81\u51fa\uff0c\u6bcf\u6b21\u5979\u89aa\u543b\u6211\u6642\uff0c\u5979\u5c31\u6703\u7528\u9019\u5169\u9846\u7259\u9f52\u6233\u6211\u7684\u81c9\u9830\u3002\u6c92\u4ec0\u9ebc\u6bd4\u88ab\u9019\u4e9b\u65e2\u89aa\u5207\u53c8\u4e0d\u597d\u770b\u7684\u89aa\u621a\u64c1\u62b1\u66f4\u4ee4\u6211\u5bb3\u6015\u7684\u4e86\u3002\u6709\u4e00\u6b21\uff0c\u6211\u88ab\u6bcd\u89aa\u62b1\u5728\u61f7\u88e1\uff0c\u5979\u5011\u554f\u6211\uff0c\u5979\u5011\u5169\u500b\u8ab0\u6bd4\u8f03\u6f02\u4eae\u3002\u6211\u4ed4\u7d30\u5730\u5be9\u8996\u4e86\u5979\u4eec\u7684\u81c9\uff0c\u7136\u5f8c\u6c89\u601d\u5730\u56de\u7b54\uff0c\u6307\u5411\u5176\u4e2d\u4e00\u500b\uff0c\u8aaa\uff1a\u300c\u9019\u500b\u6c92\u6709\u53e6\u4e00\u500b\u90a3\u9ebc\u4e11\u3002\u300d\", newText: \"\u6211\u6240\u505a\u7684\u4e00\u5207\u503c\u5f97\u7a31\u8b9a\u7684\u4e8b\uff0c\u53ea\u6703\u8b93\u6211\u7684\u7236\u6bcd\u66f4\u52a0\u60b2\u75db\u5730\u611f\u53d7\u5230\u4ed6\u5011\u7684\u640d\u5931\u3002\u56e0\u6b64\uff0c\u6211\u9577\u5927\u5f8c\u5c0d\u81ea\u5df1\u7f3a\u4e4f\u81ea\u4fe1\u3002\u4f46\u5982\u679c\u8981\u5f9e\u6211\u4ecd\u7136\u8a18\u61b6\u6df1\u523b\u7684\u4e00\u4ef6\u4e8b\u4f86\u5224\u65b7\uff0c\u6211\u7d55\u975e\u611a\u8822\u7684\u5b69\u5b50\u3002\u6709\u4e00\u5929\uff0c\u5e02\u8b70\u54e1\u5011\u7d93\u904e\u4e00\u689d\u6211\u8207\u5176\u4ed6\u7537\u5b69\u73a9\u800d\u7684\u8857\u9053\u3002\u9019\u4e9b\u5c0a\u656c\u9577\u8005\u7684\u9577\u8005\u4e2d\uff0c\u6700\u5e74\u9577\u7684\u4e00\u4f4d\u2014\u2014\u4e00\u4f4d\u5bcc\u6709\u7684\u5e02\u6c11\u2014\u2014\u505c\u4e0b\u4f86\uff0c\u7d66\u6211\u5011\u6bcf\u4eba\u4e00\u679a\u9280\u5e63\u3002\u8d70\u5230\u6211\u9762\u524d\uff0c\u4ed6\u7a81\u7136\u505c\u4e0b\uff0c\u547d\u4ee4\u9053\uff1a\u300c\u770b\u8457\u6211\u7684\u773c\u775b\u3002\u300d\u6211\u8207\u4ed6\u7684\u76ee\u5149\u76f8\u9047\uff0c\u4f38\u51fa\u624b\u53bb\u63a5\u90a3\u679a\u73cd\u8cb4\u7684\u786c\u5e63\uff0c\u537b\u767c\u73fe\u4ed6\u8a9e\u51fa\u9a5a\u4eba\uff0c\u8aaa\u9053\uff1a\u300c\u4e0d\uff0c\u5225\u592a\u591a\uff0c\u4f60\u5f9e\u6211\u9019\u5152\u5f97\u4e0d\u5230\u4ec0\u9ebc\uff0c\u4f60\u592a\u8070\u660e\u4e86\u3002\u300d\u4ed6\u5011\u5e38\u8aaa\u4e00\u500b\u95dc\u65bc\u6211\u7684\u6709\u8da3\u6545\u4e8b\u3002\u6211\u6709\u5169\u500b\u81c9\u4e0a\u6709\u76ba\u7d0b\u7684\u8001\u963f\u59e8\uff0c\u5176\u4e2d\u4e00\u4f4d\u6709\u5169\u9846\u50cf\u8c61\u7259\u4e00\u6a23\u7a81\u51fa\u7684\u7259\u9f52\uff0c\u6bcf\u6b21\u89aa\u543b\u6211\u6642\uff0c\u5979\u5c31\u6703\u628a\u7259\u9f52\u6df1\u6df1\u5730\u57cb\u5728\u6211\u7684\u81c9\u9830\u4e0a\u3002\u6c92\u6709\u4ec0\u9ebc\u6bd4\u88ab\u9019\u4e9b\u65e2\u6df1\u60c5\u53c8\u96e3\u770b\u7684\u89aa\u621a\u64c1\u62b1\u66f4\u4ee4\u6211\u5bb3\u6015\u4e86\u3002\u6709\u4e00\u6b21\uff0c\u6211\u88ab\u6bcd\u89aa\u62b1\u5728\u61f7\u88e1\uff0c\u5979\u5011\u554f\u6211\u54ea\u500b\u66f4\u6f02\u4eae\u3002\u6211\u4ed4\u7d30\u5730\u89c0\u5bdf\u4e86\u4ed6\u5011\u7684\u81c9\uff0c\u7136\u5f8c\u6c89\u601d\u5730\u6307\u8457\u5176\u4e2d\u4e00\u4f4d\uff0c\u8aaa\u9053\uff1a\u300c\u9019\u500b\u4e0d\u50cf\u53e6\u4e00\u500b\u90a3\u9ebc\u919c\u3002\u300d\" },\n  { oldText: \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\u81f3\u95dc\u91cd\u8981\u5730\u4ef0\u8cf4\u767c\u660e\uff0c\u5b83\u662f\u4eba\u985e\u5275\u9020\u6027\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\u767c\u660e\u7684\u6700\u7d42\u76ee\u7684\u5728\u65bc\u5b8c\u5168\u638c\u63e1\u5fc3\u9748\u5c0d\u7269\u8cea\u4e16\u754c\u7684\u652f\u914d\uff0c\u5c07\u81ea\u7136\u7684\u529b\u91cf\u99d5\u99ad\u65bc\u4eba\u985e\u7684\u9700\u6c42\u3002\u9019\u9805\u8271\u5de8\u4efb\u52d9\u5e38\u4f7f\u767c\u660e\u5bb6\u98fd\u53d7\u8aa4\u89e3\u8207\u5ffd\u7565\uff0c\u4f46\u4ed6\u5011\u5728\u767c\u63ee\u81ea\u8eab\u80fd\u529b\u7684\u904e\u7a0b\u4e2d\u7372\u5f97\u4e86\u5145\u6c9b\u7684\u56de\u5831\uff0c\u4e26\u4ee5\u8eab\u70ba\u300c\u7279\u6b0a\u968e\u7d1a\u300d\u7684\u4e00\u4efd\u5b50\u800c\u611f\u5230\u81ea\u8c6a\u3002\u56e0\u70ba\u6b63\u662f\u4ed6\u5011\uff0c\u8b93\u4eba\u985e\u514d\u65bc\u5728\u7121\u60c5\u81ea\u7136\u5143\u7d20\u7684\u82e6\u82e6\u9b25\u722d\u4e2d\u6ec5\u7d55\u3002\", newText: \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\uff0c\u4ef0\u8cf4\u65bc\u767c\u660e\u3002\u767c\u660e\u662f\u4eba\u985e\u5275\u9020\u529b\u7684\u7d50\u6676\uff0c\u5176\u6700\u7d42\u76ee\u6a19\u662f\u8b93\u4eba\u985e\u638c\u63e1\u7269\u8cea\u4e16\u754c\uff0c\u99d5\u99ad\u81ea\u7136\u529b\u91cf\uff0c\u6eff\u8db3\u81ea\u8eab\u9700\u6c42\u3002\u9019\u9805\u8271\u9245\u7684\u4efb\u52d9\u843d\u5728\u767c\u660e\u5bb6\u7684\u80a9\u4e0a\u3002\u4ed6\u5011\u5e38\u88ab\u8aa4\u89e3\uff0c\u4e5f\u5e38\u5f97\u4e0d\u5230\u56de\u5831\u3002\u4f46\u767c\u660e\u5bb6\u5728\u767c\u63ee\u5176\u80fd\u529b\uff0c\u610f\u8b58\u5230\u81ea\u5df1\u662f\u7279\u6b8a\u4e14\u6709\u7279\u6b0a\u7684\u4e00\u7fa4\u6642\uff0c\u5df2\u7372\u5f97\u4e86\u6700\u5927\u7684\u56de\u5831\u3002\u6c92\u6709\u4ed6\u5011\uff0c\u4eba\u985e\u7a2e\u65cf\u65e9\u5df2\u5728\u8207\u81ea\u7136\u74b0\u5883\u7684\u6297\u722d\u4e2d\u6ec5\u4ea1\u3002\" },\n  { oldText: \"\u5c31\u6211\u500b\u4eba\u800c\u8a00\uff0c\u6211\u5df2\u7d93\u5145\u5206\u9ad4\u9a57\u904e\u9019\u7a2e\u7121\u8207\u502b\u6bd4\u7684\u559c\u6085\uff0c\u751a\u81f3\u591a\u5e74\u4f86\u6211\u7684\u751f\u6d3b\u5e7e\u4e4e\u7b49\u540c\u65bc\u6301\u7e8c\u7684\u72c2\u559c\u3002\u6211\u88ab\u8a8d\u70ba\u662f\u4e00\u4f4d\u52e4\u596e\u5de5\u4f5c\u8005\uff0c\u4e5f\u8a31\u7684\u78ba\u5982\u6b64\uff0c\u56e0\u70ba\u6211\u5c07\u5e7e\u4e4e\u6240\u6709\u6e05\u9192\u7684\u6642\u9593\u90fd\u6295\u5165\u601d\u8003\u3002\u4f46\u5982\u679c\u5c07\u5de5\u4f5c\u5b9a\u7fa9\u70ba\u6839\u64da\u56b4\u683c\u898f\u5247\u5728\u7279\u5b9a\u6642\u9593\u5167\u5b8c\u6210\u7684\u660e\u78ba\u4efb\u52d9\uff0c\u90a3\u9ebc\u6211\u53ef\u80fd\u662f\u6700\u61f6\u60f0\u7684\u4eba\u3002\u4efb\u4f55\u88ab\u8feb\u7684\u52aa\u529b\u90fd\u9700\u8981\u4ed8\u51fa\u751f\u547d\u80fd\u91cf\u7684\u4ee3\u50f9\u3002\u6211\u5f9e\u672a\u4ed8\u51fa\u904e\u9019\u6a23\u7684\u4ee3\u50f9\u3002\u76f8\u53cd\u5730\uff0c\u6211\u5f9e\u601d\u60f3\u4e2d\u7372\u5f97\u4e86\u5de8\u5927\u7684\u6eff\u8db3\u3002\", newText: \"\u5c31\u500b\u4eba\u800c\u8a00\uff0c\u6211\u5df2\u6c89\u6d78\u65bc\u7cbe\u7dfb\u7684\u6a02\u8da3\uff0c\u751a\u81f3\u9054\u5230\u98fd\u548c\u7684\u7a0b\u5ea6\uff0c\u591a\u5e74\u4f86\u6211\u7684\u751f\u6d3b\u5e7e\u4e4e\u7b49\u540c\u65bc\u6301\u7e8c\u7684\u72c2\u559c\u3002\u6211\u88ab\u8a8d\u70ba\u662f\u6700\u52aa\u529b\u5de5\u4f5c\u7684\u4eba\u4e4b\u4e00\uff0c\u4e5f\u8a31\u7684\u78ba\u5982\u6b64\uff0c\u5982\u679c\u5c07\u601d\u8003\u8996\u70ba\u52de\u52d5\u7684\u8a71\uff0c\u56e0\u70ba\u6211\u5e7e\u4e4e\u5c07\u6240\u6709\u6e05\u9192\u6642\u9593\u90fd\u5949\u737b\u65bc\u6b64\u3002\u7136\u800c\uff0c\u82e5\u5c07\u5de5\u4f5c\u5b9a\u7fa9\u70ba\u5728\u898f\u5b9a\u7684\u6642\u9593\u5167\u6309\u7167\u56b4\u683c\u898f\u5247\u9032\u884c\u7684\u660e\u78ba\u8868\u73fe\uff0c\u90a3\u6211\u53ef\u80fd\u5c31\u662f\u6700\u61f6\u60f0\u7684\u4eba\u3002\u4efb\u4f55\u88ab\u8feb\u7684\u52aa\u529b\u90fd\u8981\u6c42\u72a7\u7272\u751f\u547d\u80fd\u91cf\uff0c\u6211\u5f9e\u672a\u4ed8\u51fa\u9019\u6a23\u7684\u4ee3\u50f9\u3002\u76f8\u53cd\uff0c\u6211\u5f9e\u81ea\u8eab\u601d\u60f3\u4e2d\u7372\u5f97\u6ecb\u990a\u3002\" },\n  { oldText: \"\u70ba\u4e86\u8b93\u9019\u7cfb\u5217\u6587\u7ae0\u80fd\u4ee5\u9023\u8cab\u4e14\u5fe0\u5be6\u7684\u65b9\u5f0f\u5448\u73fe\u6211\u7684\u7d93\u6b77\uff0c\u9019\u4e9b\u6587\u7ae0\u5c07\u5728\u300a\u96fb\u6c23\u5be6\u9a57\u5bb6\u300b\u7de8\u8f2f\u7684\u5354\u52a9\u4e0b\u5b8c\u6210\uff0c\u4e26\u4e3b\u8981\u9396\u5b9a\u5e74\u8f15\u7537\u6027\u8b80\u8005\u3002\u56e0\u6b64\uff0c\u6211\u4e0d\u5f97\u4e0d\u52c9\u5f37\u56de\u9867\u9752\u6625\u6642\u671f\u7684\u611f\u53d7\uff0c\u4ee5\u53ca\u90a3\u4e9b\u4fc3\u6210\u6211\u8077\u696d\u751f\u6daf\u7684\u74b0\u5883\u8207\u4e8b\u4ef6\u3002\", newText: \"\u70ba\u4e86\u8b93\u5e74\u8f15\u8b80\u8005\u5011\u80fd\u4ee5\u9023\u8cab\u4e14\u5fe0\u5be6\u7684\u65b9\u5f0f\uff0c\u85c9\u7531\u300a\u96fb\u6c23\u5be6\u9a57\u5bb6\u300b\u7de8\u8f2f\u5011\u7684\u5354\u52a9\uff0c\u4e86\u89e3\u6211\u7684\u6d3b\u52d5\u8a18\u9304\uff0c\u6211\u4e0d\u5f97\u4e0d\u56de\u9867\u81ea\u5df1\u7684\u9752\u6625\u6642\u671f\uff0c\u4ee5\u53ca\u90a3\u4e9b\u5f71\u97ff\u6211\u8077\u696d\u751f\u6daf\u7684\u74b0\u5883\u548c\u4e8b\u4ef6\u3002\" },\n  { oldText: \"\u6211\u5011\u6700\u521d\u7684\u52aa\u529b\u7d14\u7cb9\u662f\u6191\u8457\u672c\u80fd\uff0c\u4f86\u81ea\u5145\u6eff\u6d3b\u529b\u3001\u4e0d\u53d7\u62d8\u675f\u7684\u60f3\u50cf\u529b\u9a45\u4f7f\u3002\u96a8\u8457\u5e74\u9f61\u589e\u9577\uff0c\u7406\u6027\u958b\u59cb\u4f54\u64da\u4e3b\u5c0e\u5730\u4f4d\uff0c\u6211\u5011\u8b8a\u5f97\u8d8a\u4f86\u8d8a\u6709\u689d\u7406\u3001\u6709\u8a08\u756b\u6027\u3002\u4f46\u90a3\u4e9b\u65e9\u671f\u7684\u885d\u52d5\uff0c\u96d6\u7136\u4e0d\u898b\u5f97\u7acb\u5373\u958b\u82b1\u7d50\u679c\uff0c\u537b\u81f3\u95dc\u91cd\u8981\uff0c\u53ef\u80fd\u5851\u9020\u6211\u5011\u7684\u4eba\u751f\u8ecc\u8de1\u3002\u4e8b\u5be6\u4e0a\uff0c\u6211\u73fe\u5728\u89ba\u5f97\uff0c\u5982\u679c\u7576\u6642\u6211\u80fd\u7406\u89e3\u4e26\u57f9\u990a\u5b83\u5011\uff0c\u800c\u4e0d\u662f\u58d3\u6291\u5b83\u5011\uff0c\u6211\u7559\u7d66\u4e16\u754c\u7684\u907a\u7522\u5c07\u6703\u66f4\u6709\u50f9\u503c\u3002\u7136\u800c\uff0c\u76f4\u5230\u6210\u5e74\u5f8c\uff0c\u6211\u624d\u610f\u8b58\u5230\u81ea\u5df1\u662f\u4e00\u4f4d\u767c\u660e\u5bb6\u3002\", newText: \"\u6211\u5011\u6700\u521d\u7684\u52aa\u529b\u5b8c\u5168\u662f\u6191\u76f4\u89ba\uff0c\u7531\u6d3b\u6f51\u4e0d\u53d7\u62d8\u675f\u7684\u60f3\u50cf\u529b\u9a45\u4f7f\u3002\u96a8\u8457\u5e74\u9f61\u589e\u9577\uff0c\u7406\u6027\u958b\u59cb\u4f54\u64da\u4e3b\u5c0e\u5730\u4f4d\uff0c\u6211\u5011\u8b8a\u5f97\u8d8a\u4f86\u8d8a\u6709\u7cfb\u7d71\u548c\u8a2d\u8a08\u611f\u3002\u7136\u800c\uff0c\u90a3\u4e9b\u65e9\u671f\u7684\u885d\u52d5\uff0c\u5373\u4f7f\u6c92\u6709\u7acb\u5373\u7522\u751f\u6210\u679c\uff0c\u537b\u6975\u70ba\u91cd\u8981\uff0c\u5b83\u5011\u53ef\u80fd\u5851\u9020\u6211\u5011\u7684\u4eba\u751f\u9053\u8def\u3002\u4e8b\u5be6\u4e0a\uff0c\u6211\u73fe\u5728\u610f\u8b58\u5230\uff0c\u5982\u679c\u80fd\u7406\u89e3\u4e26\u57f9\u990a\u9019\u4e9b\u885d\u52d5\uff0c\u800c\u4e0d\u662f\u58d3\u6291\u5b83\u5011\uff0c\u6211\u7559\u7d66\u4e16\u754c\u7684\u907a\u7522\u5c07\u66f4\u5177\u50f9\u503c\u3002\u76f4\u5230\u6210\u5e74\u5f8c\uff0c\u6211\u624d\u610f\u8b58\u5230\u81ea\u5df1\u662f\u4e00\u4f4d\u767c\u660e\u5bb6\u3002\" },\n  { oldText: \"\u9019\u4e00\u5207\u90fd\u6709\u8a31\u591a\u539f\u56e0\u3002\u9996\u5148\uff0c\u6211\u6709\u4e00\u500b\u5929\u8ce6\u7570\u7a1f\u7684\u5f1f\u5f1f\uff0c\u4ed6\u7684\u5b58\u5728\u662f\u751f\u7269\u5b78\u7814\u7a76\u7121\u6cd5\u89e3\u91cb\u7684\u7f55\u898b\u5fc3\u667a\u73fe\u8c61\u4e4b\u4e00\u3002\u4ed6\u7684\u65e9\u901d\u4ee4\u7236\u6bcd\u60b2\u75db\u6b32\u7d55\u3002\u6211\u5011\u990a\u4e86\u4e00\u5339\u7531\u4e00\u4f4d\u646f\u53cb\u8d08\u9001\u7684\u99ac\uff0c\u9019\u5339\u963f\u62c9\u4f2f\u8840\u7d71\u7684\u99ff\u99ac\u667a\u529b\u8d85\u7fa4\uff0c\u5099\u53d7\u5168\u5bb6\u4eba\u5bf5\u611b\u8207\u5475\u8b77\uff0c\u9084\u66fe\u5728\u4e00\u6b21\u975e\u51e1\u4e8b\u4ef6\u4e2d\u6551\u4e86\u7236\u89aa\u7684\u6027\u547d\u3002\u4e00\u500b\u51ac\u591c\uff0c\u7236\u89aa\u5949\u547d\u524d\u5f80\u5c71\u5340\u57f7\u884c\u7dca\u6025\u4efb\u52d9\uff0c\u800c\u90a3\u88e1\u6709\u72fc\u51fa\u6c92\u3002\u99ac\u5339\u53d7\u9a5a\u5931\u63a7\uff0c\u5c07\u7236\u89aa\u91cd\u91cd\u6454\u5012\u5728\u5730\u3002\u7260\u56de\u5230\u5bb6\u4e2d\uff0c\u904d\u9ad4\u9c57\u50b7\uff0c\u7cbe\u75b2\u529b\u76e1\uff0c\u4f46\u8b66\u5831\u97ff\u8d77\u5f8c\uff0c\u7acb\u523b\u53c8\u8fd4\u56de\u4e8b\u6545\u73fe\u5834\u3002\u641c\u6551\u968a\u5c1a\u672a\u51fa\u767c\uff0c\u5c31\u9047\u5230\u4e86\u5df2\u6062\u5fa9\u610f\u8b58\u4e26\u91cd\u65b0\u9a0e\u4e0a\u99ac\u7684\u7236\u89aa\u3002\u4ed6\u6beb\u4e0d\u77e5\u60c5\uff0c\u81ea\u5df1\u5df2\u5728\u96ea\u5730\u88e1\u660f\u8ff7\u6578\u5c0f\u6642\u3002\u9019\u5339\u99ac\u5f8c\u4f86\u9020\u6210\u5f1f\u5f1f\u53d7\u50b7\uff0c\u6700\u7d42\u5c0e\u81f4\u4ed6\u55aa\u751f\u3002\u6211\u89aa\u773c\u76ee\u7779\u4e86\u9019\u5834\u60b2\u5287\uff0c\u5118\u7ba1\u6642\u9694\u4e94\u5341\u516d\u5e74\uff0c\u6211\u5c0d\u7576\u6642\u60c5\u666f\u7684\u8a18\u61b6\u4ecd\u6b77\u6b77\u5728\u76ee\u3002\u6bcf\u7576\u60f3\u8d77\u5f1f\u5f1f\u7684\u6210\u5c31\uff0c\u6211\u6240\u505a\u7684\u4e00\u5207\u52aa\u529b\u90fd\u986f\u5f97\u5fae\u4e0d\u8db3\u9053\u3002\", newText: \"\u9019\u4e00\u5207\u6e90\u65bc\u8a31\u591a\u539f\u56e0\u3002\u9996\u5148\uff0c\u6211\u6709\u4e00\u500b\u5929\u8ce6\u7570\u7a1f\u7684\u54e5\u54e5\uff0c\u4ed6\u64c1\u6709\u7f55\u898b\u7684\u667a\u6167\uff0c\u662f\u751f\u7269\u5b78\u7814\u7a76\u7121\u6cd5\u89e3\u91cb\u7684\u73fe\u8c61\u3002\u4ed6\u7684\u65e9\u901d\u8b93\u7236\u6bcd\u5fc3\u788e\u3002\u6211\u5011\u9084\u6709\u4e00\u5339\u7531\u646f\u53cb\u8d08\u9001\u7684\u99ac\uff0c\u5b83\u662f\u4e00\u5339\u9ad8\u8cb4\u7684\u963f\u62c9\u4f2f\u99ac\uff0c\u5e7e\u4e4e\u64c1\u6709\u4eba\u7684\u667a\u6167\uff0c\u53d7\u5230\u5168\u5bb6\u4eba\u7684\u611b\u8b77\u3002\u5b83\u66fe\u5728\u5371\u6025\u6642\u523b\u6551\u904e\u6211\u7236\u89aa\u4e00\u547d\u3002\u6709\u4e00\u6b21\uff0c\u4e00\u500b\u51ac\u591c\uff0c\u6211\u7236\u89aa\u88ab\u7dca\u6025\u53ec\u559a\u51fa\u4efb\u52d9\uff0c\u5728\u7a7f\u8d8a\u72fc\u7fa4\u51fa\u6c92\u7684\u5c71\u5340\u6642\uff0c\u99ac\u53d7\u9a5a\u800c\u9003\uff0c\u5c07\u4ed6\u62cb\u4e0b\u3002\u99ac\u56de\u5230\u5bb6\u5f8c\uff0c\u6eff\u8eab\u662f\u8840\uff0c\u7cbe\u75b2\u529b\u7aed\uff0c\u4f46\u807d\u5230\u8b66\u5831\u5f8c\uff0c\u7acb\u523b\u53c8\u885d\u51fa\u53bb\uff0c\u56de\u5230\u4e8b\u767c\u5730\u9ede\u3002\u641c\u5c0b\u968a\u51fa\u767c\u4e0d\u4e45\u4fbf\u627e\u5230\u4e86\u6211\u7236\u89aa\uff0c\u4ed6\u5df2\u7d93\u6062\u5fa9\u610f\u8b58\uff0c\u91cd\u65b0\u9a0e\u4e0a\u99ac\uff0c\u537b\u4e0d\u77e5\u9053\u81ea\u5df1\u66fe\u5728\u96ea\u5730\u88e1\u660f\u8ff7\u4e86\u5e7e\u500b\u5c0f\u6642\u3002\u9019\u5339\u99ac\u4e5f\u9020\u6210\u6211\u54e5\u54e5\u7684\u610f\u5916\u6b7b\u4ea1\u3002\u6211\u89aa\u773c\u76ee\u7779\u4e86\u9019\u5834\u60b2\u5287\uff0c\u5373\u4f7f\u904e\u4e86\u4e94\u5341\u516d\u5e74\uff0c\u6211\u4ecd\u7136\u6b77\u6b77\u5728\u76ee\u3002\u56de\u60f3\u8d77\u4ed6\u7684\u6210\u5c31\uff0c\u8b93\u6211\u6240\u6709\u7684\u52aa\u529b\u90fd\u9eef\u7136\u5931\u8272\u3002\" },\n  { oldText: \"\u6211\u6240\u6709\u7684\u512a\u9ede\uff0c\u53ea\u6703\u8b93\u7236\u6bcd\u66f4\u52a0\u601d\u5ff5\u901d\u53bb\u7684\u5b69\u5b50\u3002\u56e0\u6b64\uff0c\u6211\u6210\u9577\u904e\u7a0b\u4e2d\u59cb\u7d42\u7f3a\u4e4f\u81ea\u4fe1\u3002\u4f46\u6211\u7d55\u975e\u611a\u7b28\u7684\u7537\u5b69\uff0c\u81f3\u5c11\u5f9e\u6211\u81f3\u4eca\u8a18\u61b6\u6df1\u523b\u7684\u4e8b\u4ef6\u53ef\u4ee5\u8b49\u660e\u3002\u6709\u4e00\u5929\uff0c\u5e02\u8b70\u54e1\u5011\u7d93\u904e\u4e00\u689d\u8857\u9053\uff0c\u6211\u548c\u5176\u4ed6\u7537\u5b69\u6b63\u5728\u90a3\u88e1\u73a9\u800d\u3002\u5176\u4e2d\u4e00\u4f4d\u5bcc\u6709\u7684\u9577\u8005\u505c\u4e0b\u4f86\uff0c\u7d66\u4e86\u6211\u5011\u6bcf\u4eba\u4e00\u679a\u9280\u5e63\u3002\u4ed6\u8d70\u5230\u6211\u9762\u524d\uff0c\u7a81\u7136\u505c\u4f4f\uff0c\u547d\u4ee4\u9053\uff1a\u300c\u770b\u8457\u6211\u7684\u773c\u775b\u3002\u300d\u6211\u8207\u4ed6\u76ee\u5149\u76f8\u63a5\uff0c\u4f38\u51fa\u624b\u53bb\u63a5\u90a3\u679a\u73cd\u8cb4\u7684\u786c\u5e63\uff0c\u7136\u800c\uff0c\u4ee4\u6211\u6cae\u55aa\u7684\u662f\uff0c\u4ed6\u8aaa\uff1a\u300c\u4e0d\uff0c\u4e0d\u7528\u4e86\uff0c\u4f60\u5f9e\u6211\u9019\u88e1\u5f97\u4e0d\u5230\u4efb\u4f55\u6771\u897f\uff0c\u4f60\u592a\u8070\u660e\u4e86\u3002\u300d\u4ed6\u5011\u7d93\u5e38\u8b1b\u8ff0\u4e00\u500b\u95dc\u65bc\u6211\u7684\u6709\u8da3\u6545\u4e8b\u3002\u6211\u6709\u5169\u4f4d\u81c9\u4e0a\u6709\u76ba\u7d0b\u7684\u963f\u59e8\uff0c\u5176\u4e2d\u4e00\u4f4d\u9577\u8457\u5169\u9846\u50cf\u8c61\u7259\u4e00\u6a23\u7a81\u51fa\u7684\u7259\u9f52\uff0c\u6bcf\u6b21\u5979\u89aa\u543b\u6211\u6642\uff0c\u5979\u90fd\u6703\u7528\u9019\u5169\u9846\u7259\u9f52\u6233\u6211\u7684\u81c9\u9830\u3002\u6c92\u6709\u4ec0\u9ebc\u6bd4\u88ab\u9019\u4e9b\u65e2\u89aa\u8fd1\u53c8\u4e0d\u597d\u770b\u7684\u89aa\u621a\u64c1\u62b1\u66f4\u4ee4\u6211\u5bb3\u6015\u7684\u4e86\u3002\u6709\u4e00\u6b21\uff0c\u6211\u88ab\u6bcd\u89aa\u62b1\u5728\u61f7\u88e1\uff0c\u5979\u5011\u554f\u6211\uff0c\u5979\u5011\u5169\u500b\u8ab0\u6bd4\u8f03\u6f02\u4eae\u3002\u6211\u4ed4\u7d30\u5730\u5be9\u8996\u4e86\u5979\u4eec\u7684\u81c9\uff0c\u7136\u5f8c\u6c89\u601d\u5730\u56de\u7b54\uff0c\u6307\u5411\u5176\u4e2d\u4e00\u4f4d\uff0c\u8aaa\uff1a\u300c\u9019\u500b\u6c92\u6709\u53e6\u4e00\u500b\u90a3\u9ebc\u919c\u3002\u300d\", newText: \"\u6211\u505a\u7684\u4e00\u5207\u503c\u5f97\u7a31\u8b9a\u7684\u4e8b\uff0c\u53ea\u6703\u8b93\u7236\u6bcd\u66f4\u52a0\u60b2\u75db\u5730\u601d\u5ff5\u6211\u3002\u56e0\u6b64\uff0c\u6211\u9577\u5927\u5f8c\u7f3a\u4e4f\u81ea\u4fe1\u3002\u4f46\u6211\u4e26\u975e\u611a\u8822\u7684\u5b69\u5b50\u3002\u6709\u4e00\u6b21\uff0c\u5e02\u8b70\u54e1\u5011\u7d93\u904e\u6211\u8207\u5176\u4ed6\u7537\u5b69\u73a9\u800d\u7684\u8857\u9053\u3002\u6700\u5e74\u9577\u7684\u4e00\u4f4d\u5bcc\u6709\u7684\u5e02\u6c11\u505c\u4e0b\uff0c\u7d66\u6211\u5011\u6bcf\u4eba\u4e00\u679a\u9280\u5e63\u3002\u4ed6\u8d70\u8fd1\u6211\uff0c\u547d\u4ee4\u9053\uff1a\u300c\u770b\u8457\u6211\u7684\u773c\u775b\u3002\u300d\u6211\u8207\u4ed6\u5c0d\u8996\uff0c\u4f38\u624b\u53bb\u63a5\u9280\u5e63\uff0c\u537b\u807d\u898b\u4ed6\u8aaa\uff1a\u300c\u4e0d\uff0c\u4f60\u5f9e\u6211\u9019\u5152\u5f97\u4e0d\u5230\u4ec0\u9ebc\uff0c\u4f60\u592a\u8070\u660e\u4e86\u3002\u300d\u4ed6\u5011\u5e38\u8aaa\u4e00\u500b\u95dc\u65bc\u6211\u7684\u6709\u8da3\u6545\u4e8b\u3002\u6211\u6709\u5169\u500b\u81c9\u4e0a\u6709\u76ba\u7d0b\u7684\u8001\u963f\u59e8\uff0c\u5176\u4e2d\u4e00\u4f4d\u7259\u9f52\u50cf\u8c61\u7259\u4e00\u6a23\u7a81\u51fa\uff0c\u6bcf\u6b21\u89aa\u543b\u6211\u6642\uff0c\u5979\u90fd\u6703\u7528\u7259\u9f52\u8e6d\u6211\u7684\u81c9\u9830\u3002\u6211\u975e\u5e38\u5bb3\u6015\u5979\u5011\u7684\u64c1\u62b1\u3002\u6709\u4e00\u6b21\uff0c\u6bcd\u89aa\u62b1\u8457\u6211\uff0c\u5979\u5011\u554f\u6211\u54ea\u500b\u66f4\u6f02\u4eae\u3002\u6211\u4ed4\u7d30\u89c0\u5bdf\u4e86\u5979\u4eec\u7684\u81c9\uff0c\u6307\u8457\u5176\u4e2d\u4e00\u4f4d\u8aaa\uff1a\u300c\u9019\u500b\u4e0d\u50cf\u53e6\u4e00\u500b\u90a3\u9ebc\u919c\u3002\u300d\" },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet appliedCount = 0;\nfor (const para of paragraphs.items) {\n  const text = para.text;\n  const match = replacements.find((r) => r.oldText === text);\n  if (match) {\n    const range = para.getRange();\n    range.insertText(match.newText, Word.InsertLocation.replace);\n    appliedCount++;\n  }\n}\nawait context.sync();\n\nreturn `applied ${appliedCount} of ${replacements.length} replacements`;", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ OldText = \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\u6975\u5ea6\u4f9d\u8cf4\u767c\u660e\u3002\u9019\u662f\u4eba\u985e\u5275\u9020\u6027\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\u5176\u6700\u7d42\u76ee\u7684\u5728\u65bc\u5b8c\u5168\u638c\u63a7\u5fc3\u9748\u5c0d\u7269\u8cea\u4e16\u754c\u7684\u652f\u914d\uff0c\u5c07\u81ea\u7136\u7684\u529b\u91cf\u99d5\u99ad\u65bc\u4eba\u985e\u7684\u9700\u6c42\u3002\u9019\u662f\u767c\u660e\u5bb6\u7684\u8271\u5de8\u4efb\u52d9\uff0c\u4ed6\u5e38\u88ab\u8aa4\u89e3\u4e14\u4e0d\u53d7\u91cd\u8996\u3002\u4f46\u4ed6\u5f9e\u5176\u80fd\u529b\u7684\u6109\u5feb\u904b\u7528\u4e2d\uff0c\u4ee5\u53ca\u5f9e\u77e5\u9053\u81ea\u5df1\u5c6c\u65bc\u4e00\u500b\u6c92\u6709\u4e86\u5c31\u4eba\u985e\u65e9\u5c31\u6703\u5728\u8207\u7121\u60c5\u5143\u7d20\u7684\u82e6\u82e6\u9b25\u722d\u4e2d\u6ec5\u7d55\u7684\u6975\u5176\u7279\u6b0a\u968e\u7d1a\u4e2d\uff0c\u627e\u5230\u4e86\u5145\u8db3\u7684\u88dc\u511f\u3002\"; NewText = \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\uff0c\u81f3\u95dc\u91cd\u8981\u5730\u4f9d\u8cf4\u65bc\u767c\u660e\u3002\u5b83\u662f\u4eba\u985e\u5275\u9020\u6027\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\u5176\u6700\u7d42\u76ee\u7684\u662f\u5fc3\u9748\u5c0d\u7269\u8cea\u4e16\u754c\u7684\u5b8c\u5168\u638c\u63e1\uff0c\u5c07\u81ea\u7136\u529b\u91cf\u99d5\u99ad\u65bc\u4eba\u985e\u9700\u6c42\u3002\u9019\u662f\u767c\u660e\u5bb6\u8271\u5de8\u7684\u4efb\u52d9\uff0c\u4ed6\u7d93\u5e38\u88ab\u8aa4\u89e3\u548c\u5f97\u4e0d\u5230\u5831\u916c\u3002\u4f46\u4ed6\u537b\u5728\u6109\u6085\u5730\u904b\u7528\u5176\u80fd\u529b\u7684\u904e\u7a0b\u4e2d\uff0c\u4ee5\u53ca\u5728\u610f\u8b58\u5230\u81ea\u5df1\u662f\u90a3\u4e00\u7fa4\u7279\u6b8a\u800c\u6709\u7279\u6b0a\u7684\u968e\u7d1a\u4e2d\u7684\u4e00\u54e1\u6642\uff0c\u7372\u5f97\u4e86\u5145\u5206\u7684\u88dc\u511f\uff0c\u6c92\u6709\u4ed6\u5011\uff0c\u4eba\u985e\u7a2e\u65cf\u65e9\u5c31\u6703\u5728\u8207\u7121\u60c5\u5143\u7d20\u7684\u82e6\u82e6\u9b25\u722d\u4e2d\u6ec5\u4ea1\u3002\" },\n    @{ OldText = \"\u5c31\u6211\u500b\u4eba\u800c\u8a00\uff0c\u6211\u5df2\u7d93\u5145\u5206\u9ad4\u9a57\u904e\u9019\u7a2e\u7121\u8207\u502b\u6bd4\u7684\u6a02\u8da3\uff0c\u751a\u81f3\u591a\u5e74\u4f86\u6211\u7684\u751f\u6d3b\u5e7e\u4e4e\u7b49\u540c\u65bc\u6301\u7e8c\u7684\u72c2\u559c\u3002\u6211\u88ab\u8a8d\u70ba\u662f\u6700\u52e4\u596e\u7684\u4eba\u4e4b\u4e00\uff0c\u4e5f\u8a31\u6211\u78ba\u5be6\u5982\u6b64\uff0c\u5982\u679c\u601d\u8003\u7b49\u540c\u65bc\u52de\u52d5\uff0c\u56e0\u70ba\u6211\u5e7e\u4e4e\u628a\u6240\u6709\u6e05\u9192\u7684\u6642\u9593\u90fd\u6295\u5165\u5176\u4e2d\u3002\u4f46\u5982\u679c\u5de5\u4f5c\u88ab\u89e3\u91cb\u70ba\u6839\u64da\u56b4\u683c\u898f\u5247\u5728\u7279\u5b9a\u6642\u9593\u5167\u5b8c\u6210\u7684\u660e\u78ba\u5de5\u4f5c\uff0c\u90a3\u9ebc\u6211\u53ef\u80fd\u662f\u6700\u61f6\u60f0\u7684\u4eba\u3002\u4efb\u4f55\u88ab\u8feb\u7684\u52aa\u529b\u90fd\u9700\u8981\u72a7\u7272\u751f\u547d\u80fd\u91cf\u3002\u6211\u5f9e\u672a\u4ed8\u51fa\u904e\u9019\u6a23\u7684\u4ee3\u50f9\u3002\u76f8\u53cd\u5730\uff0c\u6211\u5f9e\u601d\u60f3\u4e2d\u7372\u76ca\u532a\u6dfa\u3002\"; NewText = \"\u5c31\u6211\u500b\u4eba\u800c\u8a00\uff0c\u6211\u5df2\u7d93\u4eab\u53d7\u5230\u9019\u7a2e\u7cbe\u7dfb\u7684\u6a02\u8da3\uff0c\u800c\u4e14\u662f\u8d85\u904e\u98fd\u548c\u7684\u7a0b\u5ea6\uff0c\u4ee5\u81f3\u65bc\u591a\u5e74\u4f86\u6211\u7684\u751f\u6d3b\u5e7e\u4e4e\u7b49\u540c\u65bc\u6301\u7e8c\u7684\u72c2\u559c\u3002\u6211\u88ab\u8a8d\u70ba\u662f\u6700\u52aa\u529b\u5de5\u4f5c\u7684\u4eba\u4e4b\u4e00\uff0c\u4e5f\u8a31\u6211\u78ba\u5be6\u662f\uff0c\u5982\u679c\u628a\u601d\u8003\u7b49\u540c\u65bc\u52de\u52d5\uff0c\u56e0\u70ba\u6211\u5e7e\u4e4e\u5c07\u6240\u6709\u7684\u6e05\u9192\u6642\u9593\u90fd\u5949\u737b\u7d66\u4e86\u5b83\u3002\u4f46\u662f\uff0c\u5982\u679c\u5de5\u4f5c\u88ab\u7406\u89e3\u70ba\u5728\u898f\u5b9a\u7684\u6642\u9593\u5167\u6309\u7167\u56b4\u683c\u7684\u898f\u5247\u9032\u884c\u7684\u660e\u78ba\u8868\u73fe\uff0c\u90a3\u9ebc\u6211\u53ef\u80fd\u662f\u6700\u61f6\u60f0\u7684\u4eba\u3002\u4efb\u4f55\u88ab\u8feb\u7684\u52aa\u529b\u90fd\u8981\u6c42\u72a7\u7272\u751f\u547d\u80fd\u91cf\u3002\u6211\u5f9e\u672a\u4ed8\u51fa\u904e\u9019\u6a23\u7684\u4ee3\u50f9\u3002\u76f8\u53cd\uff0c\u6211\u5f9e\u81ea\u5df1\u7684\u601d\u60f3\u4e2d\u7372\u5f97\u4e86\u6ecb\u990a\u3002\" },\n    @{ OldText = \"\u5728\u8a66\u5716\u4ee5\u9023\u8cab\u4e14\u5fe0\u5be6\u7684\u65b9\u5f0f\u8a18\u9304\u6211\u5728\u9019\u7cfb\u5217\u6587\u7ae0\u4e2d\u7684\u6d3b\u52d5\u6642\uff0c\u9019\u4e9b\u6587\u7ae0\u5c07\u5728\u300a\u96fb\u6c23\u5be6\u9a57\u5bb6\u300b\u7de8\u8f2f\u7684\u5354\u52a9\u4e0b\u5448\u73fe\uff0c\u4e3b\u8981\u9762\u5411\u5e74\u8f15\u7537\u6027\u8b80\u8005\uff0c\u6211\u4e0d\u5f97\u4e0d\uff0c\u5118\u7ba1\u52c9\u5f37\uff0c\u56de\u9867\u6211\u9752\u6625\u6642\u671f\u7684\u5370\u8c61\uff0c\u4ee5\u53ca\u90a3\u4e9b\u6c7a\u5b9a\u6211\u8077\u696d\u751f\u6daf\u7684\u74b0\u5883\u548c\u4e8b\u4ef6\u3002\"; NewText = \"\u5728\u5617\u8a66\u4ee5\u9023\u8cab\u4e14\u5fe0\u5be6\u7684\u65b9\u5f0f\uff0c\u85c9\u7531\u300a\u96fb\u6c23\u5be6\u9a57\u5bb6\u300b\u7de8\u8f2f\u5011\u7684\u5354\u52a9\uff0c\u5728\u672c\u7cfb\u5217\u6587\u7ae0\u4e2d\u5411\u5e74\u8f15\u8b80\u8005\u5011\u5448\u73fe\u6211\u7684\u6d3b\u52d5\u8a18\u9304\u6642\uff0c\u6211\u5fc5\u9808\u4e0d\u60c5\u9858\u5730\u56de\u9867\u6211\u7684\u9752\u5c11\u5e74\u6642\u671f\uff0c\u4ee5\u53ca\u90a3\u4e9b\u6c7a\u5b9a\u6211\u8077\u696d\u751f\u6daf\u7684\u74b0\u5883\u548c\u4e8b\u4ef6\u3002\" },\n    @{ OldText = \"\u6211\u5011\u6700\u521d\u7684\u52aa\u529b\u7d14\u7cb9\u662f\u672c\u80fd\u7684\uff0c\u4f86\u81ea\u751f\u52d5\u800c\u4e0d\u53d7\u7d04\u675f\u7684\u60f3\u50cf\u529b\u7684\u9a45\u4f7f\u3002\u96a8\u8457\u5e74\u9f61\u7684\u589e\u9577\uff0c\u7406\u6027\u958b\u59cb\u4e3b\u5f35\u81ea\u5df1\u7684\u5730\u4f4d\uff0c\u6211\u5011\u8b8a\u5f97\u8d8a\u4f86\u8d8a\u6709\u689d\u7406\u548c\u8a2d\u8a08\u6027\u3002\u4f46\u90a3\u4e9b\u65e9\u671f\u7684\u885d\u52d5\uff0c\u96d6\u7136\u4e0d\u4e00\u5b9a\u7acb\u5373\u7522\u751f\u6210\u679c\uff0c\u537b\u662f\u6700\u91cd\u8981\u7684\uff0c\u53ef\u80fd\u6703\u5851\u9020\u6211\u5011\u7684\u4eba\u751f\u8ecc\u8de1\u3002\u4e8b\u5be6\u4e0a\uff0c\u6211\u73fe\u5728\u611f\u5230\uff0c\u5982\u679c\u6211\u7576\u6642\u7406\u89e3\u4e26\u57f9\u990a\u5b83\u5011\uff0c\u800c\u4e0d\u662f\u58d3\u6291\u5b83\u5011\uff0c\u6211\u5c07\u70ba\u4e16\u754c\u7559\u4e0b\u7684\u907a\u7522\u589e\u6dfb\u5de8\u5927\u7684\u50f9\u503c\u3002\u4f46\u662f\uff0c\u76f4\u5230\u6211\u6210\u5e74\u5f8c\u624d\u610f\u8b58\u5230\u81ea\u5df1\u662f\u4e00\u4f4d\u767c\u660e\u5bb6\u3002\"; NewText = \"\u6211\u5011\u6700\u521d\u7684\u52aa\u529b\u7d14\u7cb9\u662f\u672c\u80fd\u7684\uff0c\u662f\u9bae\u6d3b\u800c\u4e0d\u53d7\u7d04\u675f\u7684\u60f3\u50cf\u529b\u7684\u9a45\u4f7f\u3002\u96a8\u8457\u6211\u5011\u5e74\u9f61\u7684\u589e\u9577\uff0c\u7406\u6027\u958b\u59cb asserted \u81ea\u6211\uff0c\u6211\u5011\u8b8a\u5f97\u8d8a\u4f86\u8d8a\u7cfb\u7d71\u5316\u548c\u6709\u8a2d\u8a08\u6027\u3002\u4f46\u90a3\u4e9b\u65e9\u671f\u7684\u885d\u52d5\uff0c\u5118\u7ba1\u6c92\u6709\u7acb\u5373\u7522\u751f\u6210\u679c\uff0c\u537b\u662f\u81f3\u95dc\u91cd\u8981\u7684\uff0c\u5b83\u5011\u53ef\u80fd\u6703\u5851\u9020\u6211\u5011\u7684\u4eba\u751f\u8ecc\u8de1\u3002\u4e8b\u5be6\u4e0a\uff0c\u6211\u73fe\u5728\u610f\u8b58\u5230\uff0c\u5982\u679c\u6211\u80fd\u7406\u89e3\u4e26\u57f9\u990a\u5b83\u5011\uff0c\u800c\u4e0d\u662f\u58d3\u6291\u5b83\u5011\uff0c\u6211\u5c07\u70ba\u4e16\u754c\u7559\u4e0b\u7684\u907a\u7522\u589e\u6dfb\u91cd\u5927\u7684\u50f9\u503c\u3002\u4f46\u76f4\u5230\u6211\u6210\u5e74\u5f8c\uff0c\u6211\u624d\u610f\u8b58\u5230\u81ea\u5df1\u662f\u4e00\u500b\u767c\u660e\u5bb6\u3002\" },\n    @{ OldText = \"\u9019\u4e00\u5207\u6e90\u65bc\u591a\u7a2e\u539f\u56e0\u3002\u9996\u5148\uff0c\u6211\u6709\u4e00\u500b\u5929\u8ce6\u7570\u7a1f\u7684\u5144\u5f1f\u2014\u2014\u90a3\u662f\u751f\u7269\u5b78\u7814\u7a76\u7121\u6cd5\u89e3\u91cb\u7684\u7f55\u898b\u5fc3\u667a\u73fe\u8c61\u4e4b\u4e00\u3002\u4ed6\u904e\u65e9\u7684\u96e2\u4e16\u8b93\u7236\u6bcd\u60b2\u75db\u6b32\u7d55\u3002\u6211\u5011\u64c1\u6709\u4e00\u5339\u7531\u4e00\u4f4d\u646f\u53cb\u8d08\u9001\u7684\u99ac\u3002\u5b83\u662f\u4e00\u5339\u963f\u62c9\u4f2f\u8840\u7d71\u7684\u96c4\u58ef\u52d5\u7269\uff0c\u64c1\u6709\u8fd1\u4e4e\u4eba\u985e\u7684\u667a\u6167\uff0c\u53d7\u5230\u5168\u5bb6\u4eba\u5bf5\u611b\u548c\u5475\u8b77\uff0c\u751a\u81f3\u66fe\u5728\u4e00\u6b21\u975e\u51e1\u7684\u4e8b\u4ef6\u4e2d\u6551\u4e86\u7236\u89aa\u7684\u6027\u547d\u3002\u4e00\u500b\u51ac\u591c\uff0c\u7236\u89aa\u63a5\u5230\u7dca\u6025\u4efb\u52d9\uff0c\u5728\u7a7f\u8d8a\u72fc\u7fa4\u51fa\u6c92\u7684\u5c71\u5340\u6642\uff0c\u99ac\u5339\u53d7\u9a5a\u9003\u8dd1\uff0c\u5c07\u4ed6\u731b\u70c8\u5730\u6454\u5012\u5728\u5730\u3002\u99ac\u5339\u56de\u5230\u5bb6\u4e2d\uff0c\u6e3e\u8eab\u662f\u8840\uff0c\u7cbe\u75b2\u529b\u76e1\uff0c\u4f46\u8b66\u5831\u97ff\u8d77\u5f8c\uff0c\u5b83\u7acb\u5373\u518d\u6b21\u51fa\u767c\uff0c\u8fd4\u56de\u539f\u8655\uff0c\u641c\u6551\u968a\u9084\u6c92\u8d70\u591a\u9060\uff0c\u5c31\u9047\u5230\u4e86\u6062\u5fa9\u610f\u8b58\u4e26\u91cd\u65b0\u9a0e\u4e0a\u99ac\u7684\u7236\u89aa\uff0c\u4ed6\u4e26\u4e0d\u77e5\u9053\u81ea\u5df1\u5df2\u7d93\u5728\u96ea\u5730\u4e2d\u8eba\u4e86\u5e7e\u500b\u5c0f\u6642\u3002\u9019\u5339\u99ac\u5c0e\u81f4\u4e86\u5f1f\u5f1f\u7684\u53d7\u50b7\uff0c\u6700\u7d42\u5c0e\u81f4\u4ed6\u6b7b\u4ea1\u3002\u6211\u89aa\u773c\u898b\u8b49\u4e86\u9019\u5834\u60b2\u5287\uff0c\u5118\u7ba1\u5df2\u7d93\u904e\u4e86\u4e94\u5341\u516d\u5e74\uff0c\u4f46\u5c0d\u90a3\u5834\u666f\u7684\u8996\u89ba\u5370\u8c61\u4ecd\u7136\u6b77\u6b77\u5728\u76ee\u3002\u56de\u61b6\u8d77\u4ed6\u7684\u6210\u5c31\uff0c\u8b93\u6211\u6240\u6709\u7684\u52aa\u529b\u90fd\u986f\u5f97\u9eef\u7136\u5931\u8272\u3002\"; NewText = \"\u9019\u7de3\u65bc\u8a31\u591a\u539f\u56e0\u3002\u9996\u5148\uff0c\u6211\u6709\u4e00\u4f4d\u5929\u8ce6\u7570\u7a1f\u7684\u54e5\u54e5\u2014\u2014\u5c6c\u65bc\u751f\u7269\u5b78\u7814\u7a76\u7121\u6cd5\u89e3\u91cb\u7684\u7f55\u898b\u5fc3\u667a\u73fe\u8c61\u4e4b\u4e00\u3002\u4ed6\u65e9\u901d\u4f7f\u7236\u6bcd\u60b2\u75db\u6b32\u7d55\u3002\u6211\u5011\u64c1\u6709\u4e00\u5339\u7531\u4e00\u4f4d\u646f\u53cb\u8d08\u9001\u7684\u99ac\u3002\u9019\u662f\u4e00\u5339\u963f\u62c9\u4f2f\u54c1\u7a2e\u7684\u96c4\u5049\u52d5\u7269\uff0c\u5e7e\u4e4e\u5177\u6709\u4eba\u985e\u7684\u667a\u6167\uff0c\u53d7\u5230\u5168\u5bb6\u4eba\u7167\u6599\u548c\u5bf5\u611b\uff0c\u66fe\u5728\u7279\u6b8a\u60c5\u6cc1\u4e0b\u6551\u904e\u6211\u7236\u89aa\u7684\u547d\u3002\u6709\u4e00\u6b21\uff0c\u5728\u4e00\u500b\u51ac\u591c\uff0c\u6211\u7236\u89aa\u88ab\u7dca\u6025\u53ec\u559a\u53bb\u57f7\u884c\u4e00\u9805\u4efb\u52d9\uff0c\u7576\u4ed6\u7a7f\u8d8a\u72fc\u7fa4\u51fa\u6c92\u7684\u5c71\u5340\u6642\uff0c\u9019\u5339\u99ac\u53d7\u9a5a\u800c\u9003\uff0c\u5c07\u4ed6\u731b\u70c8\u5730\u62cb\u5230\u5730\u4e0a\u3002\u5b83\u56de\u5230\u5bb6\u6642\uff0c\u8840\u8de1\u6591\u6591\uff0c\u7cbe\u75b2\u529b\u7aed\uff0c\u4f46\u8b66\u5831\u97ff\u8d77\u5f8c\uff0c\u5b83\u7acb\u523b\u518d\u6b21\u885d\u51fa\u53bb\uff0c\u56de\u5230\u4e8b\u767c\u5730\u9ede\u3002\u5728\u641c\u5c0b\u968a\u51fa\u767c\u4e0d\u4e45\uff0c\u4ed6\u5011\u5c31\u9047\u5230\u4e86\u6211\u7236\u89aa\uff0c\u4ed6\u5df2\u6062\u5fa9\u610f\u8b58\uff0c\u91cd\u65b0\u9a0e\u4e0a\u99ac\uff0c\u4e26\u672a\u610f\u8b58\u5230\u81ea\u5df1\u5728\u96ea\u5730\u88e1\u8eba\u4e86\u5e7e\u500b\u5c0f\u6642\u3002\u9019\u5339\u99ac\u9020\u6210\u6211\u54e5\u54e5\u7684\u50b7\u52e2\uff0c\u5c0e\u81f4\u4ed6\u6b7b\u4ea1\u3002\u6211\u89aa\u773c\u76ee\u7779\u4e86\u9019\u5834\u60b2\u5287\uff0c\u96d6\u7136\u5df2\u7d93\u904e\u4e86\u4e94\u5341\u516d\u5e74\uff0c\u4f46\u6211\u5c0d\u5b83\u7684\u8996\u89ba\u5370\u8c61\u7d72\u6beb\u672a\u6e1b\u3002\u56de\u61b6\u8d77\u4ed6\u7684\u6210\u5c31\uff0c\u8b93\u6211\u6240\u6709\u7684\u52aa\u529b\u90fd\u986f\u5f97\u9eef\u7136\u5931\u8272\u3002\" },\n    @{ OldText = \"\u6211\u6240\u505a\u7684\u4e00\u5207\u503c\u5f97\u7a31\u8b9a\u7684\u4e8b\uff0c\u53ea\u6703\u8b93\u6211\u7684\u7236\u6bcd\u66f4\u52a0\u75db\u60dc\u4ed6\u5011\u7684\u640d\u5931\u3002\u56e0\u6b64\uff0c\u6211\u5728\u6210\u9577\u904e\u7a0b\u4e2d\u7f3a\u4e4f\u81ea\u4fe1\u3002\u4f46\u6211\u7d55\u975e\u611a\u8822\u7684\u7537\u5b69\uff0c\u81f3\u5c11\u5f9e\u6211\u81f3\u4eca\u8a18\u61b6\u7336\u65b0\u7684\u4e8b\u4ef6\u4e2d\u53ef\u4ee5\u770b\u51fa\u3002\u6709\u4e00\u5929\uff0c\u5e02\u8b70\u54e1\u5011\u8def\u904e\u4e00\u689d\u8857\uff0c\u6211\u6b63\u5728\u90a3\u88e1\u548c\u5176\u4ed6\u7684\u7537\u5b69\u73a9\u800d\u3002\u9019\u4e9b\u5c0a\u8cb4\u7684\u9577\u8005\u4e2d\uff0c\u4e00\u4f4d\u5bcc\u6709\u7684\u5e02\u6c11\u505c\u4e0b\u4f86\uff0c\u7d66\u4e86\u6211\u5011\u6bcf\u500b\u4eba\u4e00\u679a\u9280\u5e63\u3002\u7576\u4ed6\u8d70\u5230\u6211\u9762\u524d\u6642\uff0c\u7a81\u7136\u505c\u4f4f\uff0c\u547d\u4ee4\u9053\uff1a\u300c\u770b\u8457\u6211\u7684\u773c\u775b\u3002\u300d\u6211\u8207\u4ed6\u7684\u76ee\u5149\u76f8\u9047\uff0c\u4f38\u51fa\u624b\u53bb\u63a5\u90a3\u679a\u73cd\u8cb4\u7684\u786c\u5e63\uff0c\u7136\u800c\uff0c\u4ee4\u6211\u6cae\u55aa\u7684\u662f\uff0c\u4ed6\u8aaa\uff1a\u300c\u4e0d\uff0c\u4e0d\u7528\u4e86\uff0c\u4f60\u5f9e\u6211\u9019\u88e1\u5f97\u4e0d\u5230\u4efb\u4f55\u6771\u897f\uff0c\u4f60\u592a\u8070\u660e\u4e86\u3002\u300d\u4ed6\u5011\u7d93\u5e38\u8b1b\u4e00\u500b\u95dc\u65bc\u6211\u7684\u6709\u8da3\u6545\u4e8b\u3002\u6211\u6709\u5169\u500b\u81c9\u4e0a\u6709\u76ba\u7d0b\u7684\u8001\u963f\u59e8\uff0c\u5176\u4e2d\u4e00\u500b\u6709\u5169\u9846\u7259\u9f52\u50cf\u8c61\u7259\u4e00\u6a23\u7a81\u51fa\uff0c\u6bcf\u6b21\u5979\u89aa\u543b\u6211\u6642\uff0c\u5979\u5c31\u6703\u7528\u9019\u5169\u9846\u7259\u9f52\u6233\u6211\u7684\u81c9\u9830\u3002\u6c92\u4ec0\u9ebc\u6bd4\u88ab\u9019\u4e9b\u65e2\u89aa\u5207\u53c8\u4e0d\u597d\u770b\u7684\u89aa\u621a\u64c1\u62b1\u66f4\u4ee4\u6211\u5bb3\u6015\u7684\u4e86\u3002\u6709\u4e00\u6b21\uff0c\u6211\u88ab\u6bcd\u89aa\u62b1\u5728\u61f7\u88e1\uff0c\u5979\u5011\u554f\u6211\uff0c\u5979\u5011\u5169\u500b\u8ab0\u6bd4\u8f03\u6f02\u4eae\u3002\u6211\u4ed4\u7d30\u5730\u5be9\u8996\u4e86\u5979\u4eec\u7684\u81c9\uff0c\u7136\u5f8c\u6c89\u601d\u5730\u56de\u7b54\uff0c\u6307\u5411\u5176\u4e2d\u4e00\u500b\uff0c\u8aaa\uff1a\u300c\u9019\u500b\u6c92\u6709\u53e6\u4e00\u500b\u90a3\u9ebc\u4e11\u3002\u300d\"; NewText = \"\u6211\u6240\u505a\u7684\u4e00\u5207\u503c\u5f97\u7a31\u8b9a\u7684\u4e8b\uff0c\u53ea\u6703\u8b93\u6211\u7684\u7236\u6bcd\u66f4\u52a0\u60b2\u75db\u5730\u611f\u53d7\u5230\u4ed6\u5011\u7684\u640d\u5931\u3002\u56e0\u6b64\uff0c\u6211\u9577\u5927\u5f8c\u5c0d\u81ea\u5df1\u7f3a\u4e4f\u81ea\u4fe1\u3002\u4f46\u5982\u679c\u8981\u5f9e\u6211\u4ecd\u7136\u8a18\u61b6\u6df1\u523b\u7684\u4e00\u4ef6\u4e8b\u4f86\u5224\u65b7\uff0c\u6211\u7d55\u975e\u611a\u8822\u7684\u5b69\u5b50\u3002\u6709\u4e00\u5929\uff0c\u5e02\u8b70\u54e1\u5011\u7d93\u904e\u4e00\u689d\u6211\u8207\u5176\u4ed6\u7537\u5b69\u73a9\u800d\u7684\u8857\u9053\u3002\u9019\u4e9b\u5c0a\u656c\u9577\u8005\u7684\u9577\u8005\u4e2d\uff0c\u6700\u5e74\u9577\u7684\u4e00\u4f4d\u2014\u2014\u4e00\u4f4d\u5bcc\u6709\u7684\u5e02\u6c11\u2014\u2014\u505c\u4e0b\u4f86\uff0c\u7d66\u6211\u5011\u6bcf\u4eba\u4e00\u679a\u9280\u5e63\u3002\u8d70\u5230\u6211\u9762\u524d\uff0c\u4ed6\u7a81\u7136\u505c\u4e0b\uff0c\u547d\u4ee4\u9053\uff1a\u300c\u770b\u8457\u6211\u7684\u773c\u775b\u3002\u300d\u6211\u8207\u4ed6\u7684\u76ee\u5149\u76f8\u9047\uff0c\u4f38\u51fa\u624b\u53bb\u63a5\u90a3\u679a\u73cd\u8cb4\u7684\u786c\u5e63\uff0c\u537b\u767c\u73fe\u4ed6\u8a9e\u51fa\u9a5a\u4eba\uff0c\u8aaa\u9053\uff1a\u300c\u4e0d\uff0c\u5225\u592a\u591a\uff0c\u4f60\u5f9e\u6211\u9019\u5152\u5f97\u4e0d\u5230\u4ec0\u9ebc\uff0c\u4f60\u592a\u8070\u660e\u4e86\u3002\u300d\u4ed6\u5011\u5e38\u8aaa\u4e00\u500b\u95dc\u65bc\u6211\u7684\u6709\u8da3\u6545\u4e8b\u3002\u6211\u6709\u5169\u500b\u81c9\u4e0a\u6709\u76ba\u7d0b\u7684\u8001\u963f\u59e8\uff0c\u5176\u4e2d\u4e00\u4f4d\u6709\u5169\u9846\u50cf\u8c61\u7259\u4e00\u6a23\u7a81\u51fa\u7684\u7259\u9f52\uff0c\u6bcf\u6b21\u89aa\u543b\u6211\u6642\uff0c\u5979\u5c31\u6703\u628a\u7259\u9f52\u6df1\u6df1\u5730\u57cb\u5728\u6211\u7684\u81c9\u9830\u4e0a\u3002\u6c92\u6709\u4ec0\u9ebc\u6bd4\u88ab\u9019\u4e9b\u65e2\u6df1\u60c5\u53c8\u96e3\u770b\u7684\u89aa\u621a\u64c1\u62b1\u66f4\u4ee4\u6211\u5bb3\u6015\u4e86\u3002\u6709\u4e00\u6b21\uff0c\u6211\u88ab\u6bcd\u89aa\u62b1\u5728\u61f7\u88e1\uff0c\u5979\u5011\u554f\u6211\u54ea\u500b\u66f4\u6f02\u4eae\u3002\u6211\u4ed4\u7d30\u5730\u89c0\u5bdf\u4e86\u4ed6\u5011\u7684\u81c9\uff0c\u7136\u5f8c\u6c89\u601d\u5730\u6307\u8457\u5176\u4e2d\u4e00\u4f4d\uff0c\u8aaa\u9053\uff1a\u300c\u9019\u500b\u4e0d\u50cf\u53e6\u4e00\u500b\u90a3\u9ebc\u919c\u3002\u300d\" },\n    @{ OldText = \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\u81f3\u95dc\u91cd\u8981\u5730\u4ef0\u8cf4\u767c\u660e\uff0c\u5b83\u662f\u4eba\u985e\u5275\u9020\u6027\u5927\u8166\u6700\u91cd\u8981\u7684\u7522\u7269\u3002\u767c\u660e\u7684\u6700\u7d42\u76ee\u7684\u5728\u65bc\u5b8c\u5168\u638c\u63e1\u5fc3\u9748\u5c0d\u7269\u8cea\u4e16\u754c\u7684\u652f\u914d\uff0c\u5c07\u81ea\u7136\u7684\u529b\u91cf\u99d5\u99ad\u65bc\u4eba\u985e\u7684\u9700\u6c42\u3002\u9019\u9805\u8271\u5de8\u4efb\u52d9\u5e38\u4f7f\u767c\u660e\u5bb6\u98fd\u53d7\u8aa4\u89e3\u8207\u5ffd\u7565\uff0c\u4f46\u4ed6\u5011\u5728\u767c\u63ee\u81ea\u8eab\u80fd\u529b\u7684\u904e\u7a0b\u4e2d\u7372\u5f97\u4e86\u5145\u6c9b\u7684\u56de\u5831\uff0c\u4e26\u4ee5\u8eab\u70ba\u300c\u7279\u6b0a\u968e\u7d1a\u300d\u7684\u4e00\u4efd\u5b50\u800c\u611f\u5230\u81ea\u8c6a\u3002\u56e0\u70ba\u6b63\u662f\u4ed6\u5011\uff0c\u8b93\u4eba\u985e\u514d\u65bc\u5728\u7121\u60c5\u81ea\u7136\u5143\u7d20\u7684\u82e6\u82e6\u9b25\u722d\u4e2d\u6ec5\u7d55\u3002\"; NewText = \"\u4eba\u985e\u7684\u9032\u6b65\u767c\u5c55\uff0c\u4ef0\u8cf4\u65bc\u767c\u660e\u3002\u767c\u660e\u662f\u4eba\u985e\u5275\u9020\u529b\u7684\u7d50\u6676\uff0c\u5176\u6700\u7d42\u76ee\u6a19\u662f\u8b93\u4eba\u985e\u638c\u63e1\u7269\u8cea\u4e16\u754c\uff0c\u99d5\u99ad\u81ea\u7136\u529b\u91cf\uff0c\u6eff\u8db3\u81ea\u8eab\u9700\u6c42\u3002\u9019\u9805\u8271\u9245\u7684\u4efb\u52d9\u843d\u5728\u767c\u660e\u5bb6\u7684\u80a9\u4e0a\u3002\u4ed6\u5011\u5e38\u88ab\u8aa4\u89e3\uff0c\u4e5f\u5e38\u5f97\u4e0d\u5230\u56de\u5831\u3002\u4f46\u767c\u660e\u5bb6\u5728\u767c\u63ee\u5176\u80fd\u529b\uff0c\u610f\u8b58\u5230\u81ea\u5df1\u662f\u7279\u6b8a\u4e14\u6709\u7279\u6b0a\u7684\u4e00\u7fa4\u6642\uff0c\u5df2\u7372\u5f97\u4e86\u6700\u5927\u7684\u56de\u5831\u3002\u6c92\u6709\u4ed6\u5011\uff0c\u4eba\u985e\u7a2e\u65cf\u65e9\u5df2\u5728\u8207\u81ea\u7136\u74b0\u5883\u7684\u6297\u722d\u4e2d\u6ec5\u4ea1\u3002\" },\n    @{ OldText = \"\u5c31\u6211\u500b\u4eba\u800c\u8a00\uff0c\u6211\u5df2\u7d93\u5145\u5206\u9ad4\u9a57\u904e\u9019\u7a2e\u7121\u8207\u502b\u6bd4\u7684\u559c\u6085\uff0c\u751a\u81f3\u591a\u5e74\u4f86\u6211\u7684\u751f\u6d3b\u5e7e\u4e4e\u7b49\u540c\u65bc\u6301\u7e8c\u7684\u72c2\u559c\u3002\u6211\u88ab\u8a8d\u70ba\u662f\u4e00\u4f4d\u52e4\u596e\u5de5\u4f5c\u8005\uff0c\u4e5f\u8a31\u7684\u78ba\u5982\u6b64\uff0c\u56e0\u70ba\u6211\u5c07\u5e7e\u4e4e\u6240\u6709\u6e05\u9192\u7684\u6642\u9593\u90fd\u6295\u5165\u601d\u8003\u3002\u4f46\u5982\u679c\u5c07\u5de5\u4f5c\u5b9a\u7fa9\u70ba\u6839\u64da\u56b4\u683c\u898f\u5247\u5728\u7279\u5b9a\u6642\u9593\u5167\u5b8c\u6210\u7684\u660e\u78ba\u4efb\u52d9\uff0c\u90a3\u9ebc\u6211\u53ef\u80fd\u662f\u6700\u61f6\u60f0\u7684\u4eba\u3002\u4efb\u4f55\u88ab\u8feb\u7684\u52aa\u529b\u90fd\u9700\u8981\u4ed8\u51fa\u751f\u547d\u80fd\u91cf\u7684\u4ee3\u50f9\u3002\u6211\u5f9e\u672a\u4ed8\u51fa\u904e\u9019\u6a23\u7684\u4ee3\u50f9\u3002\u76f8\u53cd\u5730\uff0c\u6211\u5f9e\u601d\u60f3\u4e2d\u7372\u5f97\u4e86\u5de8\u5927\u7684\u6eff\u8db3\u3002\"; NewText = \"\u5c31\u500b\u4eba\u800c\u8a00\uff0c\u6211\u5df2\u6c89\u6d78\u65bc\u7cbe\u7dfb\u7684\u6a02\u8da3\uff0c\u751a\u81f3\u9054\u5230\u98fd\u548c\u7684\u7a0b\u5ea6\uff0c\u591a\u5e74\u4f86\u6211\u7684\u751f\u6d3b\u5e7e\u4e4e\u7b49\u540c\u65bc\u6301\u7e8c\u7684\u72c2\u559c\u3002\u6211\u88ab\u8a8d\u70ba\u662f\u6700\u52aa\u529b\u5de5\u4f5c\u7684\u4eba\u4e4b\u4e00\uff0c\u4e5f\u8a31\u7684\u78ba\u5982\u6b64\uff0c\u5982\u679c\u5c07\u601d\u8003\u8996\u70ba\u52de\u52d5\u7684\u8a71\uff0c\u56e0\u70ba\u6211\u5e7e\u4e4e\u5c07\u6240\u6709\u6e05\u9192\u6642\u9593\u90fd\u5949\u737b\u65bc\u6b64\u3002\u7136\u800c\uff0c\u82e5\u5c07\u5de5\u4f5c\u5b9a\u7fa9\u70ba\u5728\u898f\u5b9a\u7684\u6642\u9593\u5167\u6309\u7167\u56b4\u683c\u898f\u5247\u9032\u884c\u7684\u660e\u78ba\u8868\u73fe\uff0c\u90a3\u6211\u53ef\u80fd\u5c31\u662f\u6700\u61f6\u60f0\u7684\u4eba\u3002\u4efb\u4f55\u88ab\u8feb\u7684\u52aa\u529b\u90fd\u8981\u6c42\u72a7\u7272\u751f\u547d\u80fd\u91cf\uff0c\u6211\u5f9e\u672a\u4ed8\u51fa\u9019\u6a23\u7684\u4ee3\u50f9\u3002\u76f8\u53cd\uff0c\u6211\u5f9e\u81ea\u8eab\u601d\u60f3\u4e2d\u7372\u5f97\u6ecb\u990a\u3002\" },\n    @{ OldText = \"\u70ba\u4e86\u8b93\u9019\u7cfb\u5217\u6587\u7ae0\u80fd\u4ee5\u9023\u8cab\u4e14\u5fe0\u5be6\u7684\u65b9\u5f0f\u5448\u73fe\u6211\u7684\u7d93\u6b77\uff0c\u9019\u4e9b\u6587\u7ae0\u5c07\u5728\u300a\u96fb\u6c23\u5be6\u9a57\u5bb6\u300b\u7de8\u8f2f\u7684\u5354\u52a9\u4e0b\u5b8c\u6210\uff0c\u4e26\u4e3b\u8981\u9396\u5b9a\u5e74\u8f15\u7537\u6027\u8b80\u8005\u3002\u56e0\u6b64\uff0c\u6211\u4e0d\u5f97\u4e0d\u52c9\u5f37\u56de\u9867\u9752\u6625\u6642\u671f\u7684\u611f\u53d7\uff0c\u4ee5\u53ca\u90a3\u4e9b\u4fc3\u6210\u6211\u8077\u696d\u751f\u6daf\u7684\u74b0\u5883\u8207\u4e8b\u4ef6\u3002\"; NewText = \"\u70ba\u4e86\u8b93\u5e74\u8f15\u8b80\u8005\u5011\u80fd\u4ee5\u9023\u8cab\u4e14\u5fe0\u5be6\u7684\u65b9\u5f0f\uff0c\u85c9\u7531\u300a\u96fb\u6c23\u5be6\u9a57\u5bb6\u300b\u7de8\u8f2f\u5011\u7684\u5354\u52a9\uff0c\u4e86\u89e3\u6211\u7684\u6d3b\u52d5\u8a18\u9304\uff0c\u6211\u4e0d\u5f97\u4e0d\u56de\u9867\u81ea\u5df1\u7684\u9752\u6625\u6642\u671f\uff0c\u4ee5\u53ca\u90a3\u4e9b\u5f71\u97ff\u6211\u8077\u696d\u751f\u6daf\u7684\u74b0\u5883\u548c\u4e8b\u4ef6\u3002\" },\n    @{ OldText = \"\u6211\u5011\u6700\u521d\u7684\u52aa\u529b\u7d14\u7cb9\u662f\u6191\u8457\u672c\u80fd\uff0c\u4f86\u81ea\u5145\u6eff\u6d3b\u529b\u3001\u4e0d\u53d7\u62d8\u675f\u7684\u60f3\u50cf\u529b\u9a45\u4f7f\u3002\u96a8\u8457\u5e74\u9f61\u589e\u9577\uff0c\u7406\u6027\u958b\u59cb\u4f54\u64da\u4e3b\u5c0e\u5730\u4f4d\uff0c\u6211\u5011\u8b8a\u5f97\u8d8a\u4f86\u8d8a\u6709\u689d\u7406\u3001\u6709\u8a08\u756b\u6027\u3002\u4f46\u90a3\u4e9b\u65e9\u671f\u7684\u885d\u52d5\uff0c\u96d6\u7136\u4e0d\u898b\u5f97\u7acb\u5373\u958b\u82b1\u7d50\u679c\uff0c\u537b\u81f3\u95dc\u91cd\u8981\uff0c\u53ef\u80fd\u5851\u9020\u6211\u5011\u7684\u4eba\u751f\u8ecc\u8de1\u3002\u4e8b\u5be6\u4e0a\uff0c\u6211\u73fe\u5728\u89ba\u5f97\uff0c\u5982\u679c\u7576\u6642\u6211\u80fd\u7406\u89e3\u4e26\u57f9\u990a\u5b83\u5011\uff0c\u800c\u4e0d\u662f\u58d3\u6291\u5b83\u5011\uff0c\u6211\u7559\u7d66\u4e16\u754c\u7684\u907a\u7522\u5c07\u6703\u66f4\u6709\u50f9\u503c\u3002\u7136\u800c\uff0c\u76f4\u5230\u6210\u5e74\u5f8c\uff0c\u6211\u624d\u610f\u8b58\u5230\u81ea\u5df1\u662f\u4e00\u4f4d\u767c\u660e\u5bb6\u3002\"; NewText = \"\u6211\u5011\u6700\u521d\u7684\u52aa\u529b\u5b8c\u5168\u662f\u6191\u76f4\u89ba\uff0c\u7531\u6d3b\u6f51\u4e0d\u53d7\u62d8\u675f\u7684\u60f3\u50cf\u529b\u9a45\u4f7f\u3002\u96a8\u8457\u5e74\u9f61\u589e\u9577\uff0c\u7406\u6027\u958b\u59cb\u4f54\u64da\u4e3b\u5c0e\u5730\u4f4d\uff0c\u6211\u5011\u8b8a\u5f97\u8d8a\u4f86\u8d8a\u6709\u7cfb\u7d71\u548c\u8a2d\u8a08\u611f\u3002\u7136\u800c\uff0c\u90a3\u4e9b\u65e9\u671f\u7684\u885d\u52d5\uff0c\u5373\u4f7f\u6c92\u6709\u7acb\u5373\u7522\u751f\u6210\u679c\uff0c\u537b\u6975\u70ba\u91cd\u8981\uff0c\u5b83\u5011\u53ef\u80fd\u5851\u9020\u6211\u5011\u7684\u4eba\u751f\u9053\u8def\u3002\u4e8b\u5be6\u4e0a\uff0c\u6211\u73fe\u5728\u610f\u8b58\u5230\uff0c\u5982\u679c\u80fd\u7406\u89e3\u4e26\u57f9\u990a\u9019\u4e9b\u885d\u52d5\uff0c\u800c\u4e0d\u662f\u58d3\u6291\u5b83\u5011\uff0c\u6211\u7559\u7d66\u4e16\u754c\u7684\u907a\u7522\u5c07\u66f4\u5177\u50f9\u503c\u3002\u76f4\u5230\u6210\u5e74\u5f8c\uff0c\u6211\u624d\u610f\u8b58\u5230\u81ea\u5df1\u662f\u4e00\u4f4d\u767c\u660e\u5bb6\u3002\" },\n    @{ OldText = \"\u9019\u4e00\u5207\u90fd\u6709\u8a31\u591a\u539f\u56e0\u3002\u9996\u5148\uff0c\u6211\u6709\u4e00\u500b\u5929\u8ce6\u7570\u7a1f\u7684\u5f1f\u5f1f\uff0c\u4ed6\u7684\u5b58\u5728\u662f\u751f\u7269\u5b78\u7814\u7a76\u7121\u6cd5\u89e3\u91cb\u7684\u7f55\u898b\u5fc3\u667a\u73fe\u8c61\u4e4b\u4e00\u3002\u4ed6\u7684\u65e9\u901d\u4ee4\u7236\u6bcd\u60b2\u75db\u6b32\u7d55\u3002\u6211\u5011\u990a\u4e86\u4e00\u5339\u7531\u4e00\u4f4d\u646f\u53cb\u8d08\u9001\u7684\u99ac\uff0c\u9019\u5339\u963f\u62c9\u4f2f\u8840\u7d71\u7684\u99ff\u99ac\u667a\u529b\u8d85\u7fa4\uff0c\u5099\u53d7\u5168\u5bb6\u4eba\u5bf5\u611b\u8207\u5475\u8b77\uff0c\u9084\u66fe\u5728\u4e00\u6b21\u975e\u51e1\u4e8b\u4ef6\u4e2d\u6551\u4e86\u7236\u89aa\u7684\u6027\u547d\u3002\u4e00\u500b\u51ac\u591c\uff0c\u7236\u89aa\u5949\u547d\u524d\u5f80\u5c71\u5340\u57f7\u884c\u7dca\u6025\u4efb\u52d9\uff0c\u800c\u90a3\u88e1\u6709\u72fc\u51fa\u6c92\u3002\u99ac\u5339\u53d7\u9a5a\u5931\u63a7\uff0c\u5c07\u7236\u89aa\u91cd\u91cd\u6454\u5012\u5728\u5730\u3002\u7260\u56de\u5230\u5bb6\u4e2d\uff0c\u904d\u9ad4\u9c57\u50b7\uff0c\u7cbe\u75b2\u529b\u76e1\uff0c\u4f46\u8b66\u5831\u97ff\u8d77\u5f8c\uff0c\u7acb\u523b\u53c8\u8fd4\u56de\u4e8b\u6545\u73fe\u5834\u3002\u641c\u6551\u968a\u5c1a\u672a\u51fa\u767c\uff0c\u5c31\u9047\u5230\u4e86\u5df2\u6062\u5fa9\u610f\u8b58\u4e26\u91cd\u65b0\u9a0e\u4e0a\u99ac\u7684\u7236\u89aa\u3002\u4ed6\u6beb\u4e0d\u77e5\u60c5\uff0c\u81ea\u5df1\u5df2\u5728\u96ea\u5730\u88e1\u660f\u8ff7\u6578\u5c0f\u6642\u3002\u9019\u5339\u99ac\u5f8c\u4f86\u9020\u6210\u5f1f\u5f1f\u53d7\u50b7\uff0c\u6700\u7d42\u5c0e\u81f4\u4ed6\u55aa\u751f\u3002\u6211\u89aa\u773c\u76ee\u7779\u4e86\u9019\u5834\u60b2\u5287\uff0c\u5118\u7ba1\u6642\u9694\u4e94\u5341\u516d\u5e74\uff0c\u6211\u5c0d\u7576\u6642\u60c5\u666f\u7684\u8a18\u61b6\u4ecd\u6b77\u6b77\u5728\u76ee\u3002\u6bcf\u7576\u60f3\u8d77\u5f1f\u5f1f\u7684\u6210\u5c31\uff0c\u6211\u6240\u505a\u7684\u4e00\u5207\u52aa\u529b\u90fd\u986f\u5f97\u5fae\u4e0d\u8db3\u9053\u3002\"; NewText = \"\u9019\u4e00\u5207\u6e90\u65bc\u8a31\u591a\u539f\u56e0\u3002\u9996\u5148\uff0c\u6211\u6709\u4e00\u500b\u5929\u8ce6\u7570\u7a1f\u7684\u54e5\u54e5\uff0c\u4ed6\u64c1\u6709\u7f55\u898b\u7684\u667a\u6167\uff0c\u662f\u751f\u7269\u5b78\u7814\u7a76\u7121\u6cd5\u89e3\u91cb\u7684\u73fe\u8c61\u3002\u4ed6\u7684\u65e9\u901d\u8b93\u7236\u6bcd\u5fc3\u788e\u3002\u6211\u5011\u9084\u6709\u4e00\u5339\u7531\u646f\u53cb\u8d08\u9001\u7684\u99ac\uff0c\u5b83\u662f\u4e00\u5339\u9ad8\u8cb4\u7684\u963f\u62c9\u4f2f\u99ac\uff0c\u5e7e\u4e4e\u64c1\u6709\u4eba\u7684\u667a\u6167\uff0c\u53d7\u5230\u5168\u5bb6\u4eba\u7684\u611b\u8b77\u3002\u5b83\u66fe\u5728\u5371\u6025\u6642\u523b\u6551\u904e\u6211\u7236\u89aa\u4e00\u547d\u3002\u6709\u4e00\u6b21\uff0c\u4e00\u500b\u51ac\u591c\uff0c\u6211\u7236\u89aa\u88ab\u7dca\u6025\u53ec\u559a\u51fa\u4efb\u52d9\uff0c\u5728\u7a7f\u8d8a\u72fc\u7fa4\u51fa\u6c92\u7684\u5c71\u5340\u6642\uff0c\u99ac\u53d7\u9a5a\u800c\u9003\uff0c\u5c07\u4ed6\u62cb\u4e0b\u3002\u99ac\u56de\u5230\u5bb6\u5f8c\uff0c\u6eff\u8eab\u662f\u8840\uff0c\u7cbe\u75b2\u529b\u7aed\uff0c\u4f46\u807d\u5230\u8b66\u5831\u5f8c\uff0c\u7acb\u523b\u53c8\u885d\u51fa\u53bb\uff0c\u56de\u5230\u4e8b\u767c\u5730\u9ede\u3002\u641c\u5c0b\u968a\u51fa\u767c\u4e0d\u4e45\u4fbf\u627e\u5230\u4e86\u6211\u7236\u89aa\uff0c\u4ed6\u5df2\u7d93\u6062\u5fa9\u610f\u8b58\uff0c\u91cd\u65b0\u9a0e\u4e0a\u99ac\uff0c\u537b\u4e0d\u77e5\u9053\u81ea\u5df1\u66fe\u5728\u96ea\u5730\u88e1\u660f\u8ff7\u4e86\u5e7e\u500b\u5c0f\u6642\u3002\u9019\u5339\u99ac\u4e5f\u9020\u6210\u6211\u54e5\u54e5\u7684\u610f\u5916\u6b7b\u4ea1\u3002\u6211\u89aa\u773c\u76ee\u7779\u4e86\u9019\u5834\u60b2\u5287\uff0c\u5373\u4f7f\u904e\u4e86\u4e94\u5341\u516d\u5e74\uff0c\u6211\u4ecd\u7136\u6b77\u6b77\u5728\u76ee\u3002\u56de\u60f3\u8d77\u4ed6\u7684\u6210\u5c31\uff0c\u8b93\u6211\u6240\u6709\u7684\u52aa\u529b\u90fd\u9eef\u7136\u5931\u8272\u3002\" },\n    @{ OldText = \"\u6211\u6240\u6709\u7684\u512a\u9ede\uff0c\u53ea\u6703\u8b93\u7236\u6bcd\u66f4\u52a0\u601d\u5ff5\u901d\u53bb\u7684\u5b69\u5b50\u3002\u56e0\u6b64\uff0c\u6211\u6210\u9577\u904e\u7a0b\u4e2d\u59cb\u7d42\u7f3a\u4e4f\u81ea\u4fe1\u3002\u4f46\u6211\u7d55\u975e\u611a\u7b28\u7684\u7537\u5b69\uff0c\u81f3\u5c11\u5f9e\u6211\u81f3\u4eca\u8a18\u61b6\u6df1\u523b\u7684\u4e8b\u4ef6\u53ef\u4ee5\u8b49\u660e\u3002\u6709\u4e00\u5929\uff0c\u5e02\u8b70\u54e1\u5011\u7d93\u904e\u4e00\u689d\u8857\u9053\uff0c\u6211\u548c\u5176\u4ed6\u7537\u5b69\u6b63\u5728\u90a3\u88e1\u73a9\u800d\u3002\u5176\u4e2d\u4e00\u4f4d\u5bcc\u6709\u7684\u9577\u8005\u505c\u4e0b\u4f86\uff0c\u7d66\u4e86\u6211\u5011\u6bcf\u4eba\u4e00\u679a\u9280\u5e63\u3002\u4ed6\u8d70\u5230\u6211\u9762\u524d\uff0c\u7a81\u7136\u505c\u4f4f\uff0c\u547d\u4ee4\u9053\uff1a\u300c\u770b\u8457\u6211\u7684\u773c\u775b\u3002\u300d\u6211\u8207\u4ed6\u76ee\u5149\u76f8\u63a5\uff0c\u4f38\u51fa\u624b\u53bb\u63a5\u90a3\u679a\u73cd\u8cb4\u7684\u786c\u5e63\uff0c\u7136\u800c\uff0c\u4ee4\u6211\u6cae\u55aa\u7684\u662f\uff0c\u4ed6\u8aaa\uff1a\u300c\u4e0d\uff0c\u4e0d\u7528\u4e86\uff0c\u4f60\u5f9e\u6211\u9019\u88e1\u5f97\u4e0d\u5230\u4efb\u4f55\u6771\u897f\uff0c\u4f60\u592a\u8070\u660e\u4e86\u3002\u300d\u4ed6\u5011\u7d93\u5e38\u8b1b\u8ff0\u4e00\u500b\u95dc\u65bc\u6211\u7684\u6709\u8da3\u6545\u4e8b\u3002\u6211\u6709\u5169\u4f4d\u81c9\u4e0a\u6709\u76ba\u7d0b\u7684\u963f\u59e8\uff0c\u5176\u4e2d\u4e00\u4f4d\u9577\u8457\u5169\u9846\u50cf\u8c61\u7259\u4e00\u6a23\u7a81\u51fa\u7684\u7259\u9f52\uff0c\u6bcf\u6b21\u5979\u89aa\u543b\u6211\u6642\uff0c\u5979\u90fd\u6703\u7528\u9019\u5169\u9846\u7259\u9f52\u6233\u6211\u7684\u81c9\u9830\u3002\u6c92\u6709\u4ec0\u9ebc\u6bd4\u88ab\u9019\u4e9b\u65e2\u89aa\u8fd1\u53c8\u4e0d\u597d\u770b\u7684\u89aa\u621a\u64c1\u62b1\u66f4\u4ee4\u6211\u5bb3\u6015\u7684\u4e86\u3002\u6709\u4e00\u6b21\uff0c\u6211\u88ab\u6bcd\u89aa\u62b1\u5728\u61f7\u88e1\uff0c\u5979\u5011\u554f\u6211\uff0c\u5979\u5011\u5169\u500b\u8ab0\u6bd4\u8f03\u6f02\u4eae\u3002\u6211\u4ed4\u7d30\u5730\u5be9\u8996\u4e86\u5979\u4eec\u7684\u81c9\uff0c\u7136\u5f8c\u6c89\u601d\u5730\u56de\u7b54\uff0c\u6307\u5411\u5176\u4e2d\u4e00\u4f4d\uff0c\u8aaa\uff1a\u300c\u9019\u500b\u6c92\u6709\u53e6\u4e00\u500b\u90a3\u9ebc\u919c\u3002\u300d\"; NewText = \"\u6211\u505a\u7684\u4e00\u5207\u503c\u5f97\u7a31\u8b9a\u7684\u4e8b\uff0c\u53ea\u6703\u8b93\u7236\u6bcd\u66f4\u52a0\u60b2\u75db\u5730\u601d\u5ff5\u6211\u3002\u56e0\u6b64\uff0c\u6211\u9577\u5927\u5f8c\u7f3a\u4e4f\u81ea\u4fe1\u3002\u4f46\u6211\u4e26\u975e\u611a\u8822\u7684\u5b69\u5b50\u3002\u6709\u4e00\u6b21\uff0c\u5e02\u8b70\u54e1\u5011\u7d93\u904e\u6211\u8207\u5176\u4ed6\u7537\u5b69\u73a9\u800d\u7684\u8857\u9053\u3002\u6700\u5e74\u9577\u7684\u4e00\u4f4d\u5bcc\u6709\u7684\u5e02\u6c11\u505c\u4e0b\uff0c\u7d66\u6211\u5011\u6bcf\u4eba\u4e00\u679a\u9280\u5e63\u3002\u4ed6\u8d70\u8fd1\u6211\uff0c\u547d\u4ee4\u9053\uff1a\u300c\u770b\u8457\u6211\u7684\u773c\u775b\u3002\u300d\u6211\u8207\u4ed6\u5c0d\u8996\uff0c\u4f38\u624b\u53bb\u63a5\u9280\u5e63\uff0c\u537b\u807d\u898b\u4ed6\u8aaa\uff1a\u300c\u4e0d\uff0c\u4f60\u5f9e\u6211\u9019\u5152\u5f97\u4e0d\u5230\u4ec0\u9ebc\uff0c\u4f60\u592a\u8070\u660e\u4e86\u3002\u300d\u4ed6\u5011\u5e38\u8aaa\u4e00\u500b\u95dc\u65bc\u6211\u7684\u6709\u8da3\u6545\u4e8b\u3002\u6211\u6709\u5169\u500b\u81c9\u4e0a\u6709\u76ba\u7d0b\u7684\u8001\u963f\u59e8\uff0c\u5176\u4e2d\u4e00\u4f4d\u7259\u9f52\u50cf\u8c61\u7259\u4e00\u6a23\u7a81\u51fa\uff0c\u6bcf\u6b21\u89aa\u543b\u6211\u6642\uff0c\u5979\u90fd\u6703\u7528\u7259\u9f52\u8e6d\u6211\u7684\u81c9\u9830\u3002\u6211\u975e\u5e38\u5bb3\u6015\u5979\u5011\u7684\u64c1\u62b1\u3002\u6709\u4e00\u6b21\uff0c\u6bcd\u89aa\u62b1\u8457\u6211\uff0c\u5979\u5011\u554f\u6211\u54ea\u500b\u66f4\u6f02\u4eae\u3002\u6211\u4ed4\u7d30\u89c0\u5bdf\u4e86\u5979\u4eec\u7684\u81c9\uff0c\u6307\u8457\u5176\u4e2d\u4e00\u4f4d\u8aaa\uff1a\u300c\u9019\u500b\u4e0d\u50cf\u53e6\u4e00\u500b\u90a3\u9ebc\u919c\u3002\u300d\" },\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.OldText\n    $find.Replacement.Text = $pair.NewText\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\nWrite-Output \"done\""}
